$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: A 46487-2023 (was row 2)
$ws.Range("A2").Value = 'A 46487-2023'
$ws.Range("B2").Value = 45197
$ws.Range("C2").Value = 46077
$ws.Range("D2").Value = 'VÄSTRA GÖTALANDS LÄN'
$ws.Range("E2").Value = 'STENUNGSUND'
$ws.Range("G2").Value = 3.4
$ws.Range("H2").Value = 18
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 9
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = 0
$ws.Range("N2").Value = 0
$ws.Range("O2").Value = 9
$ws.Range("P2").Value = 0
$ws.Range("Q2").Value = 19
$ws.Range("R2").Value = 'Duvhök`r`nEntita`r`nGulsparv`r`nKråka`r`nNordfladdermus`r`nPilgrimsfalk`r`nSpillkråka`r`nSvartvit flugsnappare`r`nTalltita`r`nStor revmossa`r`nDvärgpipistrell`r`nGråskimlig fladdermus`r`nJärnsparv`r`nKungsfågel`r`nStörre brunfladdermus`r`nStörre vattensalamander`r`nTrana`r`nTrollpipistrell`r`nVattenfladdermus'
$ws.Range("S2").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1415/artfynd/A 46487-2023 artfynd.xlsx", "A 46487-2023")'
$ws.Range("T2").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1415/kartor/A 46487-2023 karta.png", "A 46487-2023")'
$ws.Range("V2").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1415/klagomål/A 46487-2023 FSC-klagomål.docx", "A 46487-2023")'
$ws.Range("W2").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1415/klagomålsmail/A 46487-2023 FSC-klagomål mail.docx", "A 46487-2023")'
$ws.Range("X2").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1415/tillsyn/A 46487-2023 tillsynsbegäran.docx", "A 46487-2023")'
$ws.Range("Y2").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1415/tillsynsmail/A 46487-2023 tillsynsbegäran mail.docx", "A 46487-2023")'

# Row 3: A 46483-2023 (was row 3)
$ws.Range("A3").Value = 'A 46483-2023'
$ws.Range("B3").Value = 45197
$ws.Range("C3").Value = 46077
$ws.Range("D3").Value = 'VÄSTRA GÖTALANDS LÄN'
$ws.Range("E3").Value = 'STENUNGSUND'
$ws.Range("G3").Value = 15
$ws.Range("H3").Value = 9
$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 4
$ws.Range("K3").Value = 0
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = 0
$ws.Range("N3").Value = 0
$ws.Range("O3").Value = 4
$ws.Range("P3").Value = 0
$ws.Range("Q3").Value = 10
$ws.Range("R3").Value = 'Duvhök`r`nEntita`r`nJärpe`r`nSpillkråka`r`nVästlig hakmossa`r`nGrönsiska`r`nKungsfågel`r`nOrre`r`nSångsvan`r`nHuggorm'
$ws.Range("S3").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1415/artfynd/A 46483-2023 artfynd.xlsx", "A 46483-2023")'
$ws.Range("T3").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1415/kartor/A 46483-2023 karta.png", "A 46483-2023")'
$ws.Range("V3").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1415/klagomål/A 46483-2023 FSC-klagomål.docx", "A 46483-2023")'
$ws.Range("W3").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1415/klagomålsmail/A 46483-2023 FSC-klagomål mail.docx", "A 46483-2023")'
$ws.Range("X3").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1415/tillsyn/A 46483-2023 tillsynsbegäran.docx", "A 46483-2023")'
$ws.Range("Y3").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1415/tillsynsmail/A 46483-2023 tillsynsbegäran mail.docx", "A 46483-2023")'

# Row 4: A 1291-2025 (was row 4)
$ws.Range("A4").Value = 'A 1291-2025'
$ws.Range("B4").Value = 45667
$ws.Range("C4").Value = 46077
$ws.Range("D4").Value = 'VÄSTRA GÖTALANDS LÄN'
$ws.Range("E4").Value = 'STENUNGSUND'
$ws.Range("G4").Value = 9.699999999999999
$ws.Range("H4").Value = 6
$ws.Range("I4").Value = 1
$ws.Range("J4").Value = 2
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = 0
$ws.Range("N4").Value = 0
$ws.Range("O4").Value = 2
$ws.Range("P4").Value = 0
$ws.Range("Q4").Value = 7
$ws.Range("R4").Value = 'Spillkråka`r`nTalltita`r`nBjörksplintborre`r`nKungsfågel`r`nTjäder`r`nRevlummer`r`nÄkta lopplummer'
$ws.Range("S4").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1415/artfynd/A 1291-2025 artfynd.xlsx", "A 1291-2025")'
$ws.Range("T4").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1415/kartor/A 1291-2025 karta.png", "A 1291-2025")'
$ws.Range("V4").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1415/klagomål/A 1291-2025 FSC-klagomål.docx", "A 1291-2025")'
$ws.Range("W4").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1415/klagomålsmail/A 1291-2025 FSC-klagomål mail.docx", "A 1291-2025")'
$ws.Range("X4").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1415/tillsyn/A 1291-2025 tillsynsbegäran.docx", "A 1291-2025")'
$ws.Range("Y4").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1415/tillsynsmail/A 1291-2025 tillsynsbegäran mail.docx", "A 1291-2025")'

# Row 5: A 11235-2025 (was row 5)
$ws.Range("A5").Value = 'A 11235-2025'
$ws.Range("B5").Value = 45726.36305555556
$ws.Range("C5").Value = 46077
$ws.Range("D5").Value = 'VÄSTRA GÖTALANDS LÄN'
$ws.Range("E5").Value = 'STENUNGSUND'
$ws.Range("G5").Value = 4.5
$ws.Range("H5").Value = 1
$ws.Range("I5").Value = 3
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 0
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = 0
$ws.Range("N5").Value = 0
$ws.Range("O5").Value = 0
$ws.Range("P5").Value = 0
$ws.Range("Q5").Value = 4
$ws.Range("R5").Value = 'Blåmossa`r`nFällmossa`r`nKlippfrullania`r`nRevlummer'
$ws.Range("S5").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1415/artfynd/A 11235-2025 artfynd.xlsx", "A 11235-2025")'
$ws.Range("T5").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1415/kartor/A 11235-2025 karta.png", "A 11235-2025")'
$ws.Range("V5").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1415/klagomål/A 11235-2025 FSC-klagomål.docx", "A 11235-2025")'
$ws.Range("W5").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1415/klagomålsmail/A 11235-2025 FSC-klagomål mail.docx", "A 11235-2025")'
$ws.Range("X5").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1415/tillsyn/A 11235-2025 tillsynsbegäran.docx", "A 11235-2025")'
$ws.Range("Y5").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1415/tillsynsmail/A 11235-2025 tillsynsbegäran mail.docx", "A 11235-2025")'

# Row 6: A 45604-2025 (was row 6)
$ws.Range("A6").Value = 'A 45604-2025'
$ws.Range("B6").Value = 45922.69664351852
$ws.Range("C6").Value = 46077
$ws.Range("D6").Value = 'VÄSTRA GÖTALANDS LÄN'
$ws.Range("E6").Value = 'STENUNGSUND'
$ws.Range("G6").Value = 3.5
$ws.Range("H6").Value = 2
$ws.Range("I6").Value = 1
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 0
$ws.Range("L6").Value = 0
$ws.Range("M6").Value = 0
$ws.Range("N6").Value = 0
$ws.Range("O6").Value = 0
$ws.Range("P6").Value = 0
$ws.Range("Q6").Value = 3
$ws.Range("R6").Value = 'Granbräken`r`nMattlummer`r`nRevlummer'
$ws.Range("S6").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1415/artfynd/A 45604-2025 artfynd.xlsx", "A 45604-2025")'
$ws.Range("T6").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1415/kartor/A 45604-2025 karta.png", "A 45604-2025")'
$ws.Range("V6").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1415/klagomål/A 45604-2025 FSC-klagomål.docx", "A 45604-2025")'
$ws.Range("W6").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1415/klagomålsmail/A 45604-2025 FSC-klagomål mail.docx", "A 45604-2025")'
$ws.Range("X6").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1415/tillsyn/A 45604-2025 tillsynsbegäran.docx", "A 45604-2025")'
$ws.Range("Y6").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1415/tillsynsmail/A 45604-2025 tillsynsbegäran mail.docx", "A 45604-2025")'

# Row 7: A 11794-2025 (was row 7)
$ws.Range("A7").Value = 'A 11794-2025'
$ws.Range("B7").Value = 45727.73962962963
$ws.Range("C7").Value = 46077
$ws.Range("D7").Value = 'VÄSTRA GÖTALANDS LÄN'
$ws.Range("E7").Value = 'STENUNGSUND'
$ws.Range("G7").Value = 2.5
$ws.Range("H7").Value = 2
$ws.Range("I7").Value = 0
$ws.Range("J7").Value = 1
$ws.Range("K7").Value = 0
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = 0
$ws.Range("N7").Value = 0
$ws.Range("O7").Value = 1
$ws.Range("P7").Value = 0
$ws.Range("Q7").Value = 3
$ws.Range("R7").Value = 'Tallticka`r`nVanlig groda`r`nRevlummer'
$ws.Range("S7").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1415/artfynd/A 11794-2025 artfynd.xlsx", "A 11794-2025")'
$ws.Range("T7").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1415/kartor/A 11794-2025 karta.png", "A 11794-2025")'
$ws.Range("V7").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1415/klagomål/A 11794-2025 FSC-klagomål.docx", "A 11794-2025")'
$ws.Range("W7").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1415/klagomålsmail/A 11794-2025 FSC-klagomål mail.docx", "A 11794-2025")'
$ws.Range("X7").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1415/tillsyn/A 11794-2025 tillsynsbegäran.docx", "A 11794-2025")'
$ws.Range("Y7").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1415/tillsynsmail/A 11794-2025 tillsynsbegäran mail.docx", "A 11794-2025")'

# Row 8: A 44272-2024 (was row 8)
$ws.Range("A8").Value = 'A 44272-2024'
$ws.Range("B8").Value = 45573.48770833333
$ws.Range("C8").Value = 46077
$ws.Range("D8").Value = 'VÄSTRA GÖTALANDS LÄN'
$ws.Range("E8").Value = 'STENUNGSUND'
$ws.Range("G8").Value = 3.8
$ws.Range("H8").Value = 2
$ws.Range("I8").Value = 0
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 0
$ws.Range("L8").Value = 0
$ws.Range("M8").Value = 0
$ws.Range("N8").Value = 0
$ws.Range("O8").Value = 0
$ws.Range("P8").Value = 0
$ws.Range("Q8").Value = 2
$ws.Range("R8").Value = 'Mattlummer`r`nRevlummer'
$ws.Range("S8").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1415/artfynd/A 44272-2024 artfynd.xlsx", "A 44272-2024")'
$ws.Range("T8").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1415/kartor/A 44272-2024 karta.png", "A 44272-2024")'
$ws.Range("V8").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1415/klagomål/A 44272-2024 FSC-klagomål.docx", "A 44272-2024")'
$ws.Range("W8").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1415/klagomålsmail/A 44272-2024 FSC-klagomål mail.docx", "A 44272-2024")'
$ws.Range("X8").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1415/tillsyn/A 44272-2024 tillsynsbegäran.docx", "A 44272-2024")'
$ws.Range("Y8").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1415/tillsynsmail/A 44272-2024 tillsynsbegäran mail.docx", "A 44272-2024")'

# Row 9: A 60706-2023 (was row 9)
$ws.Range("A9").Value = 'A 60706-2023'
$ws.Range("B9").Value = 45260
$ws.Range("C9").Value = 46077
$ws.Range("D9").Value = 'VÄSTRA GÖTALANDS LÄN'
$ws.Range("E9").Value = 'STENUNGSUND'
$ws.Range("G9").Value = 2.4
$ws.Range("H9").Value = 1
$ws.Range("I9").Value = 1
$ws.Range("J9").Value = 1
$ws.Range("K9").Value = 0
$ws.Range("L9").Value = 0
$ws.Range("M9").Value = 0
$ws.Range("N9").Value = 0
$ws.Range("O9").Value = 1
$ws.Range("P9").Value = 0
$ws.Range("Q9").Value = 2
$ws.Range("R9").Value = 'Mindre hackspett`r`nBlåmossa'
$ws.Range("S9").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1415/artfynd/A 60706-2023 artfynd.xlsx", "A 60706-2023")'
$ws.Range("T9").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1415/kartor/A 60706-2023 karta.png", "A 60706-2023")'
$ws.Range("V9").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1415/klagomål/A 60706-2023 FSC-klagomål.docx", "A 60706-2023")'
$ws.Range("W9").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1415/klagomålsmail/A 60706-2023 FSC-klagomål mail.docx", "A 60706-2023")'
$ws.Range("X9").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1415/tillsyn/A 60706-2023 tillsynsbegäran.docx", "A 60706-2023")'
$ws.Range("Y9").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1415/tillsynsmail/A 60706-2023 tillsynsbegäran mail.docx", "A 60706-2023")'

# Row 10: A 4771-2026 (was row 10)
$ws.Range("A10").Value = 'A 4771-2026'
$ws.Range("B10").Value = 46048.54668981482
$ws.Range("C10").Value = 46077
$ws.Range("D10").Value = 'VÄSTRA GÖTALANDS LÄN'
$ws.Range("E10").Value = 'STENUNGSUND'
$ws.Range("G10").Value = 1.7
$ws.Range("H10").Value = 1
$ws.Range("I10").Value = 0
$ws.Range("J10").Value = 0
$ws.Range("K10").Value = 0
$ws.Range("L10").Value = 0
$ws.Range("M10").Value = 1
$ws.Range("N10").Value = 0
$ws.Range("O10").Value = 1
$ws.Range("P10").Value = 1
$ws.Range("Q10").Value = 2
$ws.Range("R10").Value = 'Skogsalm`r`nKungsfågel'
$ws.Range("S10").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1415/artfynd/A 4771-2026 artfynd.xlsx", "A 4771-2026")'
$ws.Range("T10").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1415/kartor/A 4771-2026 karta.png", "A 4771-2026")'
$ws.Range("V10").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1415/klagomål/A 4771-2026 FSC-klagomål.docx", "A 4771-2026")'
$ws.Range("W10").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1415/klagomålsmail/A 4771-2026 FSC-klagomål mail.docx", "A 4771-2026")'
$ws.Range("X10").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1415/tillsyn/A 4771-2026 tillsynsbegäran.docx", "A 4771-2026")'
$ws.Range("Y10").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1415/tillsynsmail/A 4771-2026 tillsynsbegäran mail.docx", "A 4771-2026")'

# Row 11: A 13912-2025 (was row 11)
$ws.Range("A11").Value = 'A 13912-2025'
$ws.Range("B11").Value = 45737
$ws.Range("C11").Value = 46077
$ws.Range("D11").Value = 'VÄSTRA GÖTALANDS LÄN'
$ws.Range("E11").Value = 'STENUNGSUND'
$ws.Range("G11").Value = 6.4
$ws.Range("H11").Value = 2
$ws.Range("I11").Value = 0
$ws.Range("J11").Value = 1
$ws.Range("K11").Value = 0
$ws.Range("L11").Value = 0
$ws.Range("M11").Value = 0
$ws.Range("N11").Value = 0
$ws.Range("O11").Value = 1
$ws.Range("P11").Value = 0
$ws.Range("Q11").Value = 2
$ws.Range("R11").Value = 'Spillkråka`r`nTjäder'
$ws.Range("S11").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1415/artfynd/A 13912-2025 artfynd.xlsx", "A 13912-2025")'
$ws.Range("T11").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1415/kartor/A 13912-2025 karta.png", "A 13912-2025")'
$ws.Range("V11").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1415/klagomål/A 13912-2025 FSC-klagomål.docx", "A 13912-2025")'
$ws.Range("W11").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1415/klagomålsmail/A 13912-2025 FSC-klagomål mail.docx", "A 13912-2025")'
$ws.Range("X11").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1415/tillsyn/A 13912-2025 tillsynsbegäran.docx", "A 13912-2025")'
$ws.Range("Y11").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1415/tillsynsmail/A 13912-2025 tillsynsbegäran mail.docx", "A 13912-2025")'

# Row 12: A 53156-2021 (was row 12)
$ws.Range("A12").Value = 'A 53156-2021'
$ws.Range("B12").Value = 44468
$ws.Range("C12").Value = 46077
$ws.Range("D12").Value = 'VÄSTRA GÖTALANDS LÄN'
$ws.Range("E12").Value = 'STENUNGSUND'
$ws.Range("G12").Value = 10.4
$ws.Range("H12").Value = 1
$ws.Range("I12").Value = 0
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 1
$ws.Range("L12").Value = 0
$ws.Range("M12").Value = 0
$ws.Range("N12").Value = 0
$ws.Range("O12").Value = 1
$ws.Range("P12").Value = 1
$ws.Range("Q12").Value = 1
$ws.Range("R12").Value = 'Knärot'
$ws.Range("S12").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1415/artfynd/A 53156-2021 artfynd.xlsx", "A 53156-2021")'
$ws.Range("T12").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1415/kartor/A 53156-2021 karta.png", "A 53156-2021")'
$ws.Range("V12").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1415/klagomål/A 53156-2021 FSC-klagomål.docx", "A 53156-2021")'
$ws.Range("W12").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1415/klagomålsmail/A 53156-2021 FSC-klagomål mail.docx", "A 53156-2021")'
$ws.Range("X12").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1415/tillsyn/A 53156-2021 tillsynsbegäran.docx", "A 53156-2021")'
$ws.Range("Y12").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1415/tillsynsmail/A 53156-2021 tillsynsbegäran mail.docx", "A 53156-2021")'

# Row 13: A 22733-2025 (was row 13)
$ws.Range("A13").Value = 'A 22733-2025'
$ws.Range("B13").Value = 45789.58555555555
$ws.Range("C13").Value = 46077
$ws.Range("D13").Value = 'VÄSTRA GÖTALANDS LÄN'
$ws.Range("E13").Value = 'STENUNGSUND'
$ws.Range("G13").Value = 1.2
$ws.Range("H13").Value = 1
$ws.Range("I13").Value = 0
$ws.Range("J13").Value = 0
$ws.Range("K13").Value = 0
$ws.Range("L13").Value = 0
$ws.Range("M13").Value = 0
$ws.Range("N13").Value = 0
$ws.Range("O13").Value = 0
$ws.Range("P13").Value = 0
$ws.Range("Q13").Value = 1
$ws.Range("R13").Value = 'Revlummer'
$ws.Range("S13").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1415/artfynd/A 22733-2025 artfynd.xlsx", "A 22733-2025")'
$ws.Range("T13").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1415/kartor/A 22733-2025 karta.png", "A 22733-2025")'
$ws.Range("V13").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1415/klagomål/A 22733-2025 FSC-klagomål.docx", "A 22733-2025")'
$ws.Range("W13").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1415/klagomålsmail/A 22733-2025 FSC-klagomål mail.docx", "A 22733-2025")'
$ws.Range("X13").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1415/tillsyn/A 22733-2025 tillsynsbegäran.docx", "A 22733-2025")'
$ws.Range("Y13").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1415/tillsynsmail/A 22733-2025 tillsynsbegäran mail.docx", "A 22733-2025")'

# Row 14: A 51807-2025 (was row 15)
$ws.Range("A14").Value = 'A 51807-2025'
$ws.Range("B14").Value = 45951.78824074074
$ws.Range("C14").Value = 46077
$ws.Range("D14").Value = 'VÄSTRA GÖTALANDS LÄN'
$ws.Range("E14").Value = 'STENUNGSUND'
$ws.Range("G14").Value = 1.2
$ws.Range("H14").Value = 0
$ws.Range("I14").Value = 1
$ws.Range("J14").Value = 0
$ws.Range("K14").Value = 0
$ws.Range("L14").Value = 0
$ws.Range("M14").Value = 0
$ws.Range("N14").Value = 0
$ws.Range("O14").Value = 0
$ws.Range("P14").Value = 0
$ws.Range("Q14").Value = 1
$ws.Range("R14").Value = 'Blomkålssvamp'
$ws.Range("S14").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1415/artfynd/A 51807-2025 artfynd.xlsx", "A 51807-2025")'
$ws.Range("T14").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1415/kartor/A 51807-2025 karta.png", "A 51807-2025")'
$ws.Range("V14").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1415/klagomål/A 51807-2025 FSC-klagomål.docx", "A 51807-2025")'
$ws.Range("W14").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1415/klagomålsmail/A 51807-2025 FSC-klagomål mail.docx", "A 51807-2025")'
$ws.Range("X14").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1415/tillsyn/A 51807-2025 tillsynsbegäran.docx", "A 51807-2025")'
$ws.Range("Y14").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1415/tillsynsmail/A 51807-2025 tillsynsbegäran mail.docx", "A 51807-2025")'

# Row 15: A 12075-2022 (was row 18)
$ws.Range("A15").Value = 'A 12075-2022'
$ws.Range("B15").Value = 44636.47275462963
$ws.Range("C15").Value = 46077
$ws.Range("D15").Value = 'VÄSTRA GÖTALANDS LÄN'
$ws.Range("E15").Value = 'STENUNGSUND'
$ws.Range("G15").Value = 7.5
$ws.Range("H15").Value = 1
$ws.Range("I15").Value = 0
$ws.Range("J15").Value = 0
$ws.Range("K15").Value = 0
$ws.Range("L15").Value = 1
$ws.Range("M15").Value = 0
$ws.Range("N15").Value = 0
$ws.Range("O15").Value = 1
$ws.Range("P15").Value = 1
$ws.Range("Q15").Value = 1
$ws.Range("R15").Value = 'Bredgentiana'
$ws.Range("S15").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1415/artfynd/A 12075-2022 artfynd.xlsx", "A 12075-2022")'
$ws.Range("T15").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1415/kartor/A 12075-2022 karta.png", "A 12075-2022")'
$ws.Range("V15").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1415/klagomål/A 12075-2022 FSC-klagomål.docx", "A 12075-2022")'
$ws.Range("W15").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1415/klagomålsmail/A 12075-2022 FSC-klagomål mail.docx", "A 12075-2022")'
$ws.Range("X15").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1415/tillsyn/A 12075-2022 tillsynsbegäran.docx", "A 12075-2022")'
$ws.Range("Y15").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1415/tillsynsmail/A 12075-2022 tillsynsbegäran mail.docx", "A 12075-2022")'

# Row 16: A 54611-2025 (was row 16)
$ws.Range("A16").Value = 'A 54611-2025'
$ws.Range("B16").Value = 45966.32216435186
$ws.Range("C16").Value = 46077
$ws.Range("D16").Value = 'VÄSTRA GÖTALANDS LÄN'
$ws.Range("E16").Value = 'STENUNGSUND'
$ws.Range("G16").Value = 10.8
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 1
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 0
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = 0
$ws.Range("N16").Value = 0
$ws.Range("O16").Value = 0
$ws.Range("P16").Value = 0
$ws.Range("Q16").Value = 1
$ws.Range("R16").Value = 'Blåmossa'
$ws.Range("S16").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1415/artfynd/A 54611-2025 artfynd.xlsx", "A 54611-2025")'
$ws.Range("T16").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1415/kartor/A 54611-2025 karta.png", "A 54611-2025")'
$ws.Range("V16").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1415/klagomål/A 54611-2025 FSC-klagomål.docx", "A 54611-2025")'
$ws.Range("W16").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1415/klagomålsmail/A 54611-2025 FSC-klagomål mail.docx", "A 54611-2025")'
$ws.Range("X16").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1415/tillsyn/A 54611-2025 tillsynsbegäran.docx", "A 54611-2025")'
$ws.Range("Y16").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1415/tillsynsmail/A 54611-2025 tillsynsbegäran mail.docx", "A 54611-2025")'

# Row 17: A 4765-2026 (was row 17)
$ws.Range("A17").Value = 'A 4765-2026'
$ws.Range("B17").Value = 46048.53950231482
$ws.Range("C17").Value = 46077
$ws.Range("D17").Value = 'VÄSTRA GÖTALANDS LÄN'
$ws.Range("E17").Value = 'STENUNGSUND'
$ws.Range("G17").Value = 1
$ws.Range("H17").Value = 1
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 0
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 1
$ws.Range("M17").Value = 0
$ws.Range("N17").Value = 0
$ws.Range("O17").Value = 1
$ws.Range("P17").Value = 1
$ws.Range("Q17").Value = 1
$ws.Range("R17").Value = 'Grönfink'
$ws.Range("S17").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1415/artfynd/A 4765-2026 artfynd.xlsx", "A 4765-2026")'
$ws.Range("T17").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1415/kartor/A 4765-2026 karta.png", "A 4765-2026")'
$ws.Range("V17").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1415/klagomål/A 4765-2026 FSC-klagomål.docx", "A 4765-2026")'
$ws.Range("W17").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1415/klagomålsmail/A 4765-2026 FSC-klagomål mail.docx", "A 4765-2026")'
$ws.Range("X17").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1415/tillsyn/A 4765-2026 tillsynsbegäran.docx", "A 4765-2026")'
$ws.Range("Y17").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1415/tillsynsmail/A 4765-2026 tillsynsbegäran mail.docx", "A 4765-2026")'

# Row 18: A 50295-2024 (was row 19)
$ws.Range("A18").Value = 'A 50295-2024'
$ws.Range("B18").Value = 45600.62445601852
$ws.Range("C18").Value = 46077
$ws.Range("D18").Value = 'VÄSTRA GÖTALANDS LÄN'
$ws.Range("E18").Value = 'STENUNGSUND'
$ws.Range("G18").Value = 2
$ws.Range("H18").Value = 1
$ws.Range("I18").Value = 0
$ws.Range("J18").Value = 0
$ws.Range("K18").Value = 0
$ws.Range("L18").Value = 0
$ws.Range("M18").Value = 0
$ws.Range("N18").Value = 0
$ws.Range("O18").Value = 0
$ws.Range("P18").Value = 0
$ws.Range("Q18").Value = 1
$ws.Range("R18").Value = 'Revlummer'
$ws.Range("S18").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1415/artfynd/A 50295-2024 artfynd.xlsx", "A 50295-2024")'
$ws.Range("T18").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1415/kartor/A 50295-2024 karta.png", "A 50295-2024")'
$ws.Range("V18").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1415/klagomål/A 50295-2024 FSC-klagomål.docx", "A 50295-2024")'
$ws.Range("W18").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1415/klagomålsmail/A 50295-2024 FSC-klagomål mail.docx", "A 50295-2024")'
$ws.Range("X18").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1415/tillsyn/A 50295-2024 tillsynsbegäran.docx", "A 50295-2024")'
$ws.Range("Y18").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1415/tillsynsmail/A 50295-2024 tillsynsbegäran mail.docx", "A 50295-2024")'

# Row 19: A 37539-2025 (was row 14)
$ws.Range("A19").Value = 'A 37539-2025'
$ws.Range("B19").Value = 45879.45538194444
$ws.Range("C19").Value = 46077
$ws.Range("D19").Value = 'VÄSTRA GÖTALANDS LÄN'
$ws.Range("E19").Value = 'STENUNGSUND'
$ws.Range("G19").Value = 1.9
$ws.Range("H19").Value = 0
$ws.Range("I19").Value = 1
$ws.Range("J19").Value = 0
$ws.Range("K19").Value = 0
$ws.Range("L19").Value = 0
$ws.Range("M19").Value = 0
$ws.Range("N19").Value = 0
$ws.Range("O19").Value = 0
$ws.Range("P19").Value = 0
$ws.Range("Q19").Value = 1
$ws.Range("R19").Value = 'Murgröna'
$ws.Range("S19").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1415/artfynd/A 37539-2025 artfynd.xlsx", "A 37539-2025")'
$ws.Range("T19").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1415/kartor/A 37539-2025 karta.png", "A 37539-2025")'
$ws.Range("V19").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1415/klagomål/A 37539-2025 FSC-klagomål.docx", "A 37539-2025")'
$ws.Range("W19").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1415/klagomålsmail/A 37539-2025 FSC-klagomål mail.docx", "A 37539-2025")'
$ws.Range("X19").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1415/tillsyn/A 37539-2025 tillsynsbegäran.docx", "A 37539-2025")'
$ws.Range("Y19").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1415/tillsynsmail/A 37539-2025 tillsynsbegäran mail.docx", "A 37539-2025")'

# Row 20: A 48215-2022 (was row 20)
$ws.Range("A20").Value = 'A 48215-2022'
$ws.Range("B20").Value = 44858
$ws.Range("C20").Value = 46077
$ws.Range("D20").Value = 'VÄSTRA GÖTALANDS LÄN'
$ws.Range("E20").Value = 'STENUNGSUND'
$ws.Range("G20").Value = 1.5
$ws.Range("H20").Value = 0
$ws.Range("I20").Value = 0
$ws.Range("J20").Value = 0
$ws.Range("K20").Value = 0
$ws.Range("L20").Value = 0
$ws.Range("M20").Value = 0
$ws.Range("N20").Value = 0
$ws.Range("O20").Value = 0
$ws.Range("P20").Value = 0
$ws.Range("Q20").Value = 0
$ws.Range("R20").Value = ''

# Row 21: A 12270-2021 (was row 21)
$ws.Range("A21").Value = 'A 12270-2021'
$ws.Range("B21").Value = 44266.91554398148
$ws.Range("C21").Value = 46077
$ws.Range("D21").Value = 'VÄSTRA GÖTALANDS LÄN'
$ws.Range("E21").Value = 'STENUNGSUND'
$ws.Range("G21").Value = 0.7
$ws.Range("H21").Value = 0
$ws.Range("I21").Value = 0
$ws.Range("J21").Value = 0
$ws.Range("K21").Value = 0
$ws.Range("L21").Value = 0
$ws.Range("M21").Value = 0
$ws.Range("N21").Value = 0
$ws.Range("O21").Value = 0
$ws.Range("P21").Value = 0
$ws.Range("Q21").Value = 0
$ws.Range("R21").Value = ''

# Row 22: A 21318-2021 (was row 22)
$ws.Range("A22").Value = 'A 21318-2021'
$ws.Range("B22").Value = 44320
$ws.Range("C22").Value = 46077
$ws.Range("D22").Value = 'VÄSTRA GÖTALANDS LÄN'
$ws.Range("E22").Value = 'STENUNGSUND'
$ws.Range("G22").Value = 1.2
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = 0
$ws.Range("N22").Value = 0
$ws.Range("O22").Value = 0
$ws.Range("P22").Value = 0
$ws.Range("Q22").Value = 0
$ws.Range("R22").Value = ''

# Row 23: A 12274-2021 (was row 23)
$ws.Range("A23").Value = 'A 12274-2021'
$ws.Range("B23").Value = 44266.91868055556
$ws.Range("C23").Value = 46077
$ws.Range("D23").Value = 'VÄSTRA GÖTALANDS LÄN'
$ws.Range("E23").Value = 'STENUNGSUND'
$ws.Range("G23").Value = 1.1
$ws.Range("H23").Value = 0
$ws.Range("I23").Value = 0
$ws.Range("J23").Value = 0
$ws.Range("K23").Value = 0
$ws.Range("L23").Value = 0
$ws.Range("M23").Value = 0
$ws.Range("N23").Value = 0
$ws.Range("O23").Value = 0
$ws.Range("P23").Value = 0
$ws.Range("Q23").Value = 0
$ws.Range("R23").Value = ''

# Row 24: A 51116-2021 (was row 24)
$ws.Range("A24").Value = 'A 51116-2021'
$ws.Range("B24").Value = 44461
$ws.Range("C24").Value = 46077
$ws.Range("D24").Value = 'VÄSTRA GÖTALANDS LÄN'
$ws.Range("E24").Value = 'STENUNGSUND'
$ws.Range("G24").Value = 2.3
$ws.Range("H24").Value = 0
$ws.Range("I24").Value = 0
$ws.Range("J24").Value = 0
$ws.Range("K24").Value = 0
$ws.Range("L24").Value = 0
$ws.Range("M24").Value = 0
$ws.Range("N24").Value = 0
$ws.Range("O24").Value = 0
$ws.Range("P24").Value = 0
$ws.Range("Q24").Value = 0
$ws.Range("R24").Value = ''

# Row 25: A 7787-2022 (was row 25)
$ws.Range("A25").Value = 'A 7787-2022'
$ws.Range("B25").Value = 44608
$ws.Range("C25").Value = 46077
$ws.Range("D25").Value = 'VÄSTRA GÖTALANDS LÄN'
$ws.Range("E25").Value = 'STENUNGSUND'
$ws.Range("G25").Value = 1.7
$ws.Range("H25").Value = 0
$ws.Range("I25").Value = 0
$ws.Range("J25").Value = 0
$ws.Range("K25").Value = 0
$ws.Range("L25").Value = 0
$ws.Range("M25").Value = 0
$ws.Range("N25").Value = 0
$ws.Range("O25").Value = 0
$ws.Range("P25").Value = 0
$ws.Range("Q25").Value = 0
$ws.Range("R25").Value = ''

# Row 26: A 41312-2021 (was row 26)
$ws.Range("A26").Value = 'A 41312-2021'
$ws.Range("B26").Value = 44424.57387731481
$ws.Range("C26").Value = 46077
$ws.Range("D26").Value = 'VÄSTRA GÖTALANDS LÄN'
$ws.Range("E26").Value = 'STENUNGSUND'
$ws.Range("G26").Value = 1.5
$ws.Range("H26").Value = 0
$ws.Range("I26").Value = 0
$ws.Range("J26").Value = 0
$ws.Range("K26").Value = 0
$ws.Range("L26").Value = 0
$ws.Range("M26").Value = 0
$ws.Range("N26").Value = 0
$ws.Range("O26").Value = 0
$ws.Range("P26").Value = 0
$ws.Range("Q26").Value = 0
$ws.Range("R26").Value = ''

# Row 27: A 3244-2022 (was row 27)
$ws.Range("A27").Value = 'A 3244-2022'
$ws.Range("B27").Value = 44582
$ws.Range("C27").Value = 46077
$ws.Range("D27").Value = 'VÄSTRA GÖTALANDS LÄN'
$ws.Range("E27").Value = 'STENUNGSUND'
$ws.Range("G27").Value = 0.8
$ws.Range("H27").Value = 0
$ws.Range("I27").Value = 0
$ws.Range("J27").Value = 0
$ws.Range("K27").Value = 0
$ws.Range("L27").Value = 0
$ws.Range("M27").Value = 0
$ws.Range("N27").Value = 0
$ws.Range("O27").Value = 0
$ws.Range("P27").Value = 0
$ws.Range("Q27").Value = 0
$ws.Range("R27").Value = ''

# Row 28: A 35927-2021 (was row 28)
$ws.Range("A28").Value = 'A 35927-2021'
$ws.Range("B28").Value = 44388.52203703704
$ws.Range("C28").Value = 46077
$ws.Range("D28").Value = 'VÄSTRA GÖTALANDS LÄN'
$ws.Range("E28").Value = 'STENUNGSUND'
$ws.Range("G28").Value = 0.9
$ws.Range("H28").Value = 0
$ws.Range("I28").Value = 0
$ws.Range("J28").Value = 0
$ws.Range("K28").Value = 0
$ws.Range("L28").Value = 0
$ws.Range("M28").Value = 0
$ws.Range("N28").Value = 0
$ws.Range("O28").Value = 0
$ws.Range("P28").Value = 0
$ws.Range("Q28").Value = 0
$ws.Range("R28").Value = ''

# Row 29: A 54601-2021 (was row 29)
$ws.Range("A29").Value = 'A 54601-2021'
$ws.Range("B29").Value = 44473
$ws.Range("C29").Value = 46077
$ws.Range("D29").Value = 'VÄSTRA GÖTALANDS LÄN'
$ws.Range("E29").Value = 'STENUNGSUND'
$ws.Range("G29").Value = 1.4
$ws.Range("H29").Value = 0
$ws.Range("I29").Value = 0
$ws.Range("J29").Value = 0
$ws.Range("K29").Value = 0
$ws.Range("L29").Value = 0
$ws.Range("M29").Value = 0
$ws.Range("N29").Value = 0
$ws.Range("O29").Value = 0
$ws.Range("P29").Value = 0
$ws.Range("Q29").Value = 0
$ws.Range("R29").Value = ''

# Row 30: A 14782-2021 (was row 30)
$ws.Range("A30").Value = 'A 14782-2021'
$ws.Range("B30").Value = 44280
$ws.Range("C30").Value = 46077
$ws.Range("D30").Value = 'VÄSTRA GÖTALANDS LÄN'
$ws.Range("E30").Value = 'STENUNGSUND'
$ws.Range("G30").Value = 5
$ws.Range("H30").Value = 0
$ws.Range("I30").Value = 0
$ws.Range("J30").Value = 0
$ws.Range("K30").Value = 0
$ws.Range("L30").Value = 0
$ws.Range("M30").Value = 0
$ws.Range("N30").Value = 0
$ws.Range("O30").Value = 0
$ws.Range("P30").Value = 0
$ws.Range("Q30").Value = 0
$ws.Range("R30").Value = ''

# Row 31: A 9273-2022 (was row 31)
$ws.Range("A31").Value = 'A 9273-2022'
$ws.Range("B31").Value = 44616
$ws.Range("C31").Value = 46077
$ws.Range("D31").Value = 'VÄSTRA GÖTALANDS LÄN'
$ws.Range("E31").Value = 'STENUNGSUND'
$ws.Range("G31").Value = 3.9
$ws.Range("H31").Value = 0
$ws.Range("I31").Value = 0
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 0
$ws.Range("L31").Value = 0
$ws.Range("M31").Value = 0
$ws.Range("N31").Value = 0
$ws.Range("O31").Value = 0
$ws.Range("P31").Value = 0
$ws.Range("Q31").Value = 0
$ws.Range("R31").Value = ''

# Row 32: A 62007-2021 (was row 33)
$ws.Range("A32").Value = 'A 62007-2021'
$ws.Range("B32").Value = 44502
$ws.Range("C32").Value = 46077
$ws.Range("D32").Value = 'VÄSTRA GÖTALANDS LÄN'
$ws.Range("E32").Value = 'STENUNGSUND'
$ws.Range("G32").Value = 1
$ws.Range("H32").Value = 0
$ws.Range("I32").Value = 0
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 0
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = 0
$ws.Range("N32").Value = 0
$ws.Range("O32").Value = 0
$ws.Range("P32").Value = 0
$ws.Range("Q32").Value = 0
$ws.Range("R32").Value = ''

# Row 33: A 6293-2022 (was row 32)
$ws.Range("A33").Value = 'A 6293-2022'
$ws.Range("B33").Value = 44600
$ws.Range("C33").Value = 46077
$ws.Range("D33").Value = 'VÄSTRA GÖTALANDS LÄN'
$ws.Range("E33").Value = 'STENUNGSUND'
$ws.Range("G33").Value = 0.3
$ws.Range("H33").Value = 0
$ws.Range("I33").Value = 0
$ws.Range("J33").Value = 0
$ws.Range("K33").Value = 0
$ws.Range("L33").Value = 0
$ws.Range("M33").Value = 0
$ws.Range("N33").Value = 0
$ws.Range("O33").Value = 0
$ws.Range("P33").Value = 0
$ws.Range("Q33").Value = 0
$ws.Range("R33").Value = ''

# Row 34: A 62220-2023 (was row 39)
$ws.Range("A34").Value = 'A 62220-2023'
$ws.Range("B34").Value = 45267.53284722222
$ws.Range("C34").Value = 46077
$ws.Range("D34").Value = 'VÄSTRA GÖTALANDS LÄN'
$ws.Range("E34").Value = 'STENUNGSUND'
$ws.Range("G34").Value = 7.6
$ws.Range("H34").Value = 0
$ws.Range("I34").Value = 0
$ws.Range("J34").Value = 0
$ws.Range("K34").Value = 0
$ws.Range("L34").Value = 0
$ws.Range("M34").Value = 0
$ws.Range("N34").Value = 0
$ws.Range("O34").Value = 0
$ws.Range("P34").Value = 0
$ws.Range("Q34").Value = 0
$ws.Range("R34").Value = ''

# Row 35: A 60158-2024 (was row 97)
$ws.Range("A35").Value = 'A 60158-2024'
$ws.Range("B35").Value = 45642.62672453704
$ws.Range("C35").Value = 46077
$ws.Range("D35").Value = 'VÄSTRA GÖTALANDS LÄN'
$ws.Range("E35").Value = 'STENUNGSUND'
$ws.Range("G35").Value = 4.4
$ws.Range("H35").Value = 0
$ws.Range("I35").Value = 0
$ws.Range("J35").Value = 0
$ws.Range("K35").Value = 0
$ws.Range("L35").Value = 0
$ws.Range("M35").Value = 0
$ws.Range("N35").Value = 0
$ws.Range("O35").Value = 0
$ws.Range("P35").Value = 0
$ws.Range("Q35").Value = 0
$ws.Range("R35").Value = ''

# Row 36: A 37866-2021 (was row 38)
$ws.Range("A36").Value = 'A 37866-2021'
$ws.Range("B36").Value = 44403.4725462963
$ws.Range("C36").Value = 46077
$ws.Range("D36").Value = 'VÄSTRA GÖTALANDS LÄN'
$ws.Range("E36").Value = 'STENUNGSUND'
$ws.Range("G36").Value = 1.5
$ws.Range("H36").Value = 0
$ws.Range("I36").Value = 0
$ws.Range("J36").Value = 0
$ws.Range("K36").Value = 0
$ws.Range("L36").Value = 0
$ws.Range("M36").Value = 0
$ws.Range("N36").Value = 0
$ws.Range("O36").Value = 0
$ws.Range("P36").Value = 0
$ws.Range("Q36").Value = 0
$ws.Range("R36").Value = ''

# Row 37: A 19059-2023 (was row 112)
$ws.Range("A37").Value = 'A 19059-2023'
$ws.Range("B37").Value = 45048
$ws.Range("C37").Value = 46077
$ws.Range("D37").Value = 'VÄSTRA GÖTALANDS LÄN'
$ws.Range("E37").Value = 'STENUNGSUND'
$ws.Range("G37").Value = 0.7
$ws.Range("H37").Value = 0
$ws.Range("I37").Value = 0
$ws.Range("J37").Value = 0
$ws.Range("K37").Value = 0
$ws.Range("L37").Value = 0
$ws.Range("M37").Value = 0
$ws.Range("N37").Value = 0
$ws.Range("O37").Value = 0
$ws.Range("P37").Value = 0
$ws.Range("Q37").Value = 0
$ws.Range("R37").Value = ''

# Row 38: A 33011-2022 (was row 54)
$ws.Range("A38").Value = 'A 33011-2022'
$ws.Range("B38").Value = 44785.31972222222
$ws.Range("C38").Value = 46077
$ws.Range("D38").Value = 'VÄSTRA GÖTALANDS LÄN'
$ws.Range("E38").Value = 'STENUNGSUND'
$ws.Range("G38").Value = 2.6
$ws.Range("H38").Value = 0
$ws.Range("I38").Value = 0
$ws.Range("J38").Value = 0
$ws.Range("K38").Value = 0
$ws.Range("L38").Value = 0
$ws.Range("M38").Value = 0
$ws.Range("N38").Value = 0
$ws.Range("O38").Value = 0
$ws.Range("P38").Value = 0
$ws.Range("Q38").Value = 0
$ws.Range("R38").Value = ''

# Row 39: A 3570-2023 (was row 46)
$ws.Range("A39").Value = 'A 3570-2023'
$ws.Range("B39").Value = 44950.49107638889
$ws.Range("C39").Value = 46077
$ws.Range("D39").Value = 'VÄSTRA GÖTALANDS LÄN'
$ws.Range("E39").Value = 'STENUNGSUND'
$ws.Range("G39").Value = 4.1
$ws.Range("H39").Value = 0
$ws.Range("I39").Value = 0
$ws.Range("J39").Value = 0
$ws.Range("K39").Value = 0
$ws.Range("L39").Value = 0
$ws.Range("M39").Value = 0
$ws.Range("N39").Value = 0
$ws.Range("O39").Value = 0
$ws.Range("P39").Value = 0
$ws.Range("Q39").Value = 0
$ws.Range("R39").Value = ''

# Row 40: A 10472-2024 (was row 80)
$ws.Range("A40").Value = 'A 10472-2024'
$ws.Range("B40").Value = 45366.39716435185
$ws.Range("C40").Value = 46077
$ws.Range("D40").Value = 'VÄSTRA GÖTALANDS LÄN'
$ws.Range("E40").Value = 'STENUNGSUND'
$ws.Range("G40").Value = 6.3
$ws.Range("H40").Value = 0
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 0
$ws.Range("L40").Value = 0
$ws.Range("M40").Value = 0
$ws.Range("N40").Value = 0
$ws.Range("O40").Value = 0
$ws.Range("P40").Value = 0
$ws.Range("Q40").Value = 0
$ws.Range("R40").Value = ''

# Row 41: A 27030-2023 (was row 37)
$ws.Range("A41").Value = 'A 27030-2023'
$ws.Range("B41").Value = 45095
$ws.Range("C41").Value = 46077
$ws.Range("D41").Value = 'VÄSTRA GÖTALANDS LÄN'
$ws.Range("E41").Value = 'STENUNGSUND'
$ws.Range("G41").Value = 3.1
$ws.Range("H41").Value = 0
$ws.Range("I41").Value = 0
$ws.Range("J41").Value = 0
$ws.Range("K41").Value = 0
$ws.Range("L41").Value = 0
$ws.Range("M41").Value = 0
$ws.Range("N41").Value = 0
$ws.Range("O41").Value = 0
$ws.Range("P41").Value = 0
$ws.Range("Q41").Value = 0
$ws.Range("R41").Value = ''

# Row 42: A 48895-2023 (was row 96)
$ws.Range("A42").Value = 'A 48895-2023'
$ws.Range("B42").Value = 45209.49538194444
$ws.Range("C42").Value = 46077
$ws.Range("D42").Value = 'VÄSTRA GÖTALANDS LÄN'
$ws.Range("E42").Value = 'STENUNGSUND'
$ws.Range("G42").Value = 3.1
$ws.Range("H42").Value = 0
$ws.Range("I42").Value = 0
$ws.Range("J42").Value = 0
$ws.Range("K42").Value = 0
$ws.Range("L42").Value = 0
$ws.Range("M42").Value = 0
$ws.Range("N42").Value = 0
$ws.Range("O42").Value = 0
$ws.Range("P42").Value = 0
$ws.Range("Q42").Value = 0
$ws.Range("R42").Value = ''

# Row 43: A 30175-2022 (was row 58)
$ws.Range("A43").Value = 'A 30175-2022'
$ws.Range("B43").Value = 44757
$ws.Range("C43").Value = 46077
$ws.Range("D43").Value = 'VÄSTRA GÖTALANDS LÄN'
$ws.Range("E43").Value = 'STENUNGSUND'
$ws.Range("G43").Value = 2.2
$ws.Range("H43").Value = 0
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 0
$ws.Range("K43").Value = 0
$ws.Range("L43").Value = 0
$ws.Range("M43").Value = 0
$ws.Range("N43").Value = 0
$ws.Range("O43").Value = 0
$ws.Range("P43").Value = 0
$ws.Range("Q43").Value = 0
$ws.Range("R43").Value = ''

# Row 44: A 23762-2024 (was row 41)
$ws.Range("A44").Value = 'A 23762-2024'
$ws.Range("B44").Value = 45455.36590277778
$ws.Range("C44").Value = 46077
$ws.Range("D44").Value = 'VÄSTRA GÖTALANDS LÄN'
$ws.Range("E44").Value = 'STENUNGSUND'
$ws.Range("G44").Value = 2
$ws.Range("H44").Value = 0
$ws.Range("I44").Value = 0
$ws.Range("J44").Value = 0
$ws.Range("K44").Value = 0
$ws.Range("L44").Value = 0
$ws.Range("M44").Value = 0
$ws.Range("N44").Value = 0
$ws.Range("O44").Value = 0
$ws.Range("P44").Value = 0
$ws.Range("Q44").Value = 0
$ws.Range("R44").Value = ''

# Row 45: A 32960-2023 (was row 59)
$ws.Range("A45").Value = 'A 32960-2023'
$ws.Range("B45").Value = 45125
$ws.Range("C45").Value = 46077
$ws.Range("D45").Value = 'VÄSTRA GÖTALANDS LÄN'
$ws.Range("E45").Value = 'STENUNGSUND'
$ws.Range("G45").Value = 4.7
$ws.Range("H45").Value = 0
$ws.Range("I45").Value = 0
$ws.Range("J45").Value = 0
$ws.Range("K45").Value = 0
$ws.Range("L45").Value = 0
$ws.Range("M45").Value = 0
$ws.Range("N45").Value = 0
$ws.Range("O45").Value = 0
$ws.Range("P45").Value = 0
$ws.Range("Q45").Value = 0
$ws.Range("R45").Value = ''

# Row 46: A 52952-2021 (was row 40)
$ws.Range("A46").Value = 'A 52952-2021'
$ws.Range("B46").Value = 44467.5972337963
$ws.Range("C46").Value = 46077
$ws.Range("D46").Value = 'VÄSTRA GÖTALANDS LÄN'
$ws.Range("E46").Value = 'STENUNGSUND'
$ws.Range("G46").Value = 0.4
$ws.Range("H46").Value = 0
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 0
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 0
$ws.Range("M46").Value = 0
$ws.Range("N46").Value = 0
$ws.Range("O46").Value = 0
$ws.Range("P46").Value = 0
$ws.Range("Q46").Value = 0
$ws.Range("R46").Value = ''

# Row 47: A 8576-2023 (was row 44)
$ws.Range("A47").Value = 'A 8576-2023'
$ws.Range("B47").Value = 44977
$ws.Range("C47").Value = 46077
$ws.Range("D47").Value = 'VÄSTRA GÖTALANDS LÄN'
$ws.Range("E47").Value = 'STENUNGSUND'
$ws.Range("G47").Value = 2.1
$ws.Range("H47").Value = 0
$ws.Range("I47").Value = 0
$ws.Range("J47").Value = 0
$ws.Range("K47").Value = 0
$ws.Range("L47").Value = 0
$ws.Range("M47").Value = 0
$ws.Range("N47").Value = 0
$ws.Range("O47").Value = 0
$ws.Range("P47").Value = 0
$ws.Range("Q47").Value = 0
$ws.Range("R47").Value = ''

# Row 48: A 19302-2025 (was row 99)
$ws.Range("A48").Value = 'A 19302-2025'
$ws.Range("B48").Value = 45769.56519675926
$ws.Range("C48").Value = 46077
$ws.Range("D48").Value = 'VÄSTRA GÖTALANDS LÄN'
$ws.Range("E48").Value = 'STENUNGSUND'
$ws.Range("G48").Value = 1.2
$ws.Range("H48").Value = 0
$ws.Range("I48").Value = 0
$ws.Range("J48").Value = 0
$ws.Range("K48").Value = 0
$ws.Range("L48").Value = 0
$ws.Range("M48").Value = 0
$ws.Range("N48").Value = 0
$ws.Range("O48").Value = 0
$ws.Range("P48").Value = 0
$ws.Range("Q48").Value = 0
$ws.Range("R48").Value = ''

# Row 49: A 33010-2022 (was row 53)
$ws.Range("A49").Value = 'A 33010-2022'
$ws.Range("B49").Value = 44785
$ws.Range("C49").Value = 46077
$ws.Range("D49").Value = 'VÄSTRA GÖTALANDS LÄN'
$ws.Range("E49").Value = 'STENUNGSUND'
$ws.Range("G49").Value = 3.4
$ws.Range("H49").Value = 0
$ws.Range("I49").Value = 0
$ws.Range("J49").Value = 0
$ws.Range("K49").Value = 0
$ws.Range("L49").Value = 0
$ws.Range("M49").Value = 0
$ws.Range("N49").Value = 0
$ws.Range("O49").Value = 0
$ws.Range("P49").Value = 0
$ws.Range("Q49").Value = 0
$ws.Range("R49").Value = ''

# Row 50: A 6913-2024 (was row 75)
$ws.Range("A50").Value = 'A 6913-2024'
$ws.Range("B50").Value = 45343.28206018519
$ws.Range("C50").Value = 46077
$ws.Range("D50").Value = 'VÄSTRA GÖTALANDS LÄN'
$ws.Range("E50").Value = 'STENUNGSUND'
$ws.Range("G50").Value = 2
$ws.Range("H50").Value = 0
$ws.Range("I50").Value = 0
$ws.Range("J50").Value = 0
$ws.Range("K50").Value = 0
$ws.Range("L50").Value = 0
$ws.Range("M50").Value = 0
$ws.Range("N50").Value = 0
$ws.Range("O50").Value = 0
$ws.Range("P50").Value = 0
$ws.Range("Q50").Value = 0
$ws.Range("R50").Value = ''

# Row 51: A 6699-2022 (was row 83)
$ws.Range("A51").Value = 'A 6699-2022'
$ws.Range("B51").Value = 44602
$ws.Range("C51").Value = 46077
$ws.Range("D51").Value = 'VÄSTRA GÖTALANDS LÄN'
$ws.Range("E51").Value = 'STENUNGSUND'
$ws.Range("G51").Value = 10
$ws.Range("H51").Value = 0
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 0
$ws.Range("K51").Value = 0
$ws.Range("L51").Value = 0
$ws.Range("M51").Value = 0
$ws.Range("N51").Value = 0
$ws.Range("O51").Value = 0
$ws.Range("P51").Value = 0
$ws.Range("Q51").Value = 0
$ws.Range("R51").Value = ''

# Row 52: A 72944-2021 (was row 48)
$ws.Range("A52").Value = 'A 72944-2021'
$ws.Range("B52").Value = 44547.87140046297
$ws.Range("C52").Value = 46077
$ws.Range("D52").Value = 'VÄSTRA GÖTALANDS LÄN'
$ws.Range("E52").Value = 'STENUNGSUND'
$ws.Range("G52").Value = 4.2
$ws.Range("H52").Value = 0
$ws.Range("I52").Value = 0
$ws.Range("J52").Value = 0
$ws.Range("K52").Value = 0
$ws.Range("L52").Value = 0
$ws.Range("M52").Value = 0
$ws.Range("N52").Value = 0
$ws.Range("O52").Value = 0
$ws.Range("P52").Value = 0
$ws.Range("Q52").Value = 0
$ws.Range("R52").Value = ''

# Row 53: A 19423-2025 (was row 116)
$ws.Range("A53").Value = 'A 19423-2025'
$ws.Range("B53").Value = 45769.76158564815
$ws.Range("C53").Value = 46077
$ws.Range("D53").Value = 'VÄSTRA GÖTALANDS LÄN'
$ws.Range("E53").Value = 'STENUNGSUND'
$ws.Range("G53").Value = 1.2
$ws.Range("H53").Value = 0
$ws.Range("I53").Value = 0
$ws.Range("J53").Value = 0
$ws.Range("K53").Value = 0
$ws.Range("L53").Value = 0
$ws.Range("M53").Value = 0
$ws.Range("N53").Value = 0
$ws.Range("O53").Value = 0
$ws.Range("P53").Value = 0
$ws.Range("Q53").Value = 0
$ws.Range("R53").Value = ''

# Row 54: A 19424-2025 (was row 117)
$ws.Range("A54").Value = 'A 19424-2025'
$ws.Range("B54").Value = 45769.7637962963
$ws.Range("C54").Value = 46077
$ws.Range("D54").Value = 'VÄSTRA GÖTALANDS LÄN'
$ws.Range("E54").Value = 'STENUNGSUND'
$ws.Range("G54").Value = 0.7
$ws.Range("H54").Value = 0
$ws.Range("I54").Value = 0
$ws.Range("J54").Value = 0
$ws.Range("K54").Value = 0
$ws.Range("L54").Value = 0
$ws.Range("M54").Value = 0
$ws.Range("N54").Value = 0
$ws.Range("O54").Value = 0
$ws.Range("P54").Value = 0
$ws.Range("Q54").Value = 0
$ws.Range("R54").Value = ''

# Row 55: A 15535-2024 (was row 109)
$ws.Range("A55").Value = 'A 15535-2024'
$ws.Range("B55").Value = 45401.610625
$ws.Range("C55").Value = 46077
$ws.Range("D55").Value = 'VÄSTRA GÖTALANDS LÄN'
$ws.Range("E55").Value = 'STENUNGSUND'
$ws.Range("G55").Value = 1.8
$ws.Range("H55").Value = 0
$ws.Range("I55").Value = 0
$ws.Range("J55").Value = 0
$ws.Range("K55").Value = 0
$ws.Range("L55").Value = 0
$ws.Range("M55").Value = 0
$ws.Range("N55").Value = 0
$ws.Range("O55").Value = 0
$ws.Range("P55").Value = 0
$ws.Range("Q55").Value = 0
$ws.Range("R55").Value = ''

# Row 56: A 3243-2022 (was row 108)
$ws.Range("A56").Value = 'A 3243-2022'
$ws.Range("B56").Value = 44582
$ws.Range("C56").Value = 46077
$ws.Range("D56").Value = 'VÄSTRA GÖTALANDS LÄN'
$ws.Range("E56").Value = 'STENUNGSUND'
$ws.Range("G56").Value = 7.6
$ws.Range("H56").Value = 0
$ws.Range("I56").Value = 0
$ws.Range("J56").Value = 0
$ws.Range("K56").Value = 0
$ws.Range("L56").Value = 0
$ws.Range("M56").Value = 0
$ws.Range("N56").Value = 0
$ws.Range("O56").Value = 0
$ws.Range("P56").Value = 0
$ws.Range("Q56").Value = 0
$ws.Range("R56").Value = ''

# Row 57: A 64992-2023 (was row 91)
$ws.Range("A57").Value = 'A 64992-2023'
$ws.Range("B57").Value = 45286
$ws.Range("C57").Value = 46077
$ws.Range("D57").Value = 'VÄSTRA GÖTALANDS LÄN'
$ws.Range("E57").Value = 'STENUNGSUND'
$ws.Range("G57").Value = 1.8
$ws.Range("H57").Value = 0
$ws.Range("I57").Value = 0
$ws.Range("J57").Value = 0
$ws.Range("K57").Value = 0
$ws.Range("L57").Value = 0
$ws.Range("M57").Value = 0
$ws.Range("N57").Value = 0
$ws.Range("O57").Value = 0
$ws.Range("P57").Value = 0
$ws.Range("Q57").Value = 0
$ws.Range("R57").Value = ''

# Row 58: A 4734-2023 (was row 105)
$ws.Range("A58").Value = 'A 4734-2023'
$ws.Range("B58").Value = 44957
$ws.Range("C58").Value = 46077
$ws.Range("D58").Value = 'VÄSTRA GÖTALANDS LÄN'
$ws.Range("E58").Value = 'STENUNGSUND'
$ws.Range("G58").Value = 1.1
$ws.Range("H58").Value = 0
$ws.Range("I58").Value = 0
$ws.Range("J58").Value = 0
$ws.Range("K58").Value = 0
$ws.Range("L58").Value = 0
$ws.Range("M58").Value = 0
$ws.Range("N58").Value = 0
$ws.Range("O58").Value = 0
$ws.Range("P58").Value = 0
$ws.Range("Q58").Value = 0
$ws.Range("R58").Value = ''

# Row 59: A 34105-2024 (was row 113)
$ws.Range("A59").Value = 'A 34105-2024'
$ws.Range("B59").Value = 45523.83888888889
$ws.Range("C59").Value = 46077
$ws.Range("D59").Value = 'VÄSTRA GÖTALANDS LÄN'
$ws.Range("E59").Value = 'STENUNGSUND'
$ws.Range("G59").Value = 2.3
$ws.Range("H59").Value = 0
$ws.Range("I59").Value = 0
$ws.Range("J59").Value = 0
$ws.Range("K59").Value = 0
$ws.Range("L59").Value = 0
$ws.Range("M59").Value = 0
$ws.Range("N59").Value = 0
$ws.Range("O59").Value = 0
$ws.Range("P59").Value = 0
$ws.Range("Q59").Value = 0
$ws.Range("R59").Value = ''

# Row 60: A 23070-2022 (was row 35)
$ws.Range("A60").Value = 'A 23070-2022'
$ws.Range("B60").Value = 44719.45615740741
$ws.Range("C60").Value = 46077
$ws.Range("D60").Value = 'VÄSTRA GÖTALANDS LÄN'
$ws.Range("E60").Value = 'STENUNGSUND'
$ws.Range("G60").Value = 4.1
$ws.Range("H60").Value = 0
$ws.Range("I60").Value = 0
$ws.Range("J60").Value = 0
$ws.Range("K60").Value = 0
$ws.Range("L60").Value = 0
$ws.Range("M60").Value = 0
$ws.Range("N60").Value = 0
$ws.Range("O60").Value = 0
$ws.Range("P60").Value = 0
$ws.Range("Q60").Value = 0
$ws.Range("R60").Value = ''

# Row 61: A 34555-2024 (was row 118)
$ws.Range("A61").Value = 'A 34555-2024'
$ws.Range("B61").Value = 45525.9071875
$ws.Range("C61").Value = 46077
$ws.Range("D61").Value = 'VÄSTRA GÖTALANDS LÄN'
$ws.Range("E61").Value = 'STENUNGSUND'
$ws.Range("G61").Value = 3.1
$ws.Range("H61").Value = 0
$ws.Range("I61").Value = 0
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 0
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = 0
$ws.Range("N61").Value = 0
$ws.Range("O61").Value = 0
$ws.Range("P61").Value = 0
$ws.Range("Q61").Value = 0
$ws.Range("R61").Value = ''

# Row 62: A 10843-2025 (was row 70)
$ws.Range("A62").Value = 'A 10843-2025'
$ws.Range("B62").Value = 45722.5778125
$ws.Range("C62").Value = 46077
$ws.Range("D62").Value = 'VÄSTRA GÖTALANDS LÄN'
$ws.Range("E62").Value = 'STENUNGSUND'
$ws.Range("G62").Value = 2.8
$ws.Range("H62").Value = 0
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("M62").Value = 0
$ws.Range("N62").Value = 0
$ws.Range("O62").Value = 0
$ws.Range("P62").Value = 0
$ws.Range("Q62").Value = 0
$ws.Range("R62").Value = ''

# Row 63: A 10849-2025 (was row 71)
$ws.Range("A63").Value = 'A 10849-2025'
$ws.Range("B63").Value = 45722.5817824074
$ws.Range("C63").Value = 46077
$ws.Range("D63").Value = 'VÄSTRA GÖTALANDS LÄN'
$ws.Range("E63").Value = 'STENUNGSUND'
$ws.Range("G63").Value = 2.9
$ws.Range("H63").Value = 0
$ws.Range("I63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("K63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("M63").Value = 0
$ws.Range("N63").Value = 0
$ws.Range("O63").Value = 0
$ws.Range("P63").Value = 0
$ws.Range("Q63").Value = 0
$ws.Range("R63").Value = ''

# Row 64: A 4578-2023 (was row 76)
$ws.Range("A64").Value = 'A 4578-2023'
$ws.Range("B64").Value = 44951
$ws.Range("C64").Value = 46077
$ws.Range("D64").Value = 'VÄSTRA GÖTALANDS LÄN'
$ws.Range("E64").Value = 'STENUNGSUND'
$ws.Range("G64").Value = 0.7
$ws.Range("H64").Value = 0
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("K64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("M64").Value = 0
$ws.Range("N64").Value = 0
$ws.Range("O64").Value = 0
$ws.Range("P64").Value = 0
$ws.Range("Q64").Value = 0
$ws.Range("R64").Value = ''

# Row 65: A 38996-2025 (was row 51)
$ws.Range("A65").Value = 'A 38996-2025'
$ws.Range("B65").Value = 45888.33400462963
$ws.Range("C65").Value = 46077
$ws.Range("D65").Value = 'VÄSTRA GÖTALANDS LÄN'
$ws.Range("E65").Value = 'STENUNGSUND'
$ws.Range("G65").Value = 2.3
$ws.Range("H65").Value = 0
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("M65").Value = 0
$ws.Range("N65").Value = 0
$ws.Range("O65").Value = 0
$ws.Range("P65").Value = 0
$ws.Range("Q65").Value = 0
$ws.Range("R65").Value = ''

# Row 66: A 20974-2025 (was row 49)
$ws.Range("A66").Value = 'A 20974-2025'
$ws.Range("B66").Value = 45777.4587962963
$ws.Range("C66").Value = 46077
$ws.Range("D66").Value = 'VÄSTRA GÖTALANDS LÄN'
$ws.Range("E66").Value = 'STENUNGSUND'
$ws.Range("G66").Value = 10.3
$ws.Range("H66").Value = 0
$ws.Range("I66").Value = 0
$ws.Range("J66").Value = 0
$ws.Range("K66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("M66").Value = 0
$ws.Range("N66").Value = 0
$ws.Range("O66").Value = 0
$ws.Range("P66").Value = 0
$ws.Range("Q66").Value = 0
$ws.Range("R66").Value = ''

# Row 67: A 14990-2025 (was row 50)
$ws.Range("A67").Value = 'A 14990-2025'
$ws.Range("B67").Value = 45743.60700231481
$ws.Range("C67").Value = 46077
$ws.Range("D67").Value = 'VÄSTRA GÖTALANDS LÄN'
$ws.Range("E67").Value = 'STENUNGSUND'
$ws.Range("G67").Value = 1.9
$ws.Range("H67").Value = 0
$ws.Range("I67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("K67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("M67").Value = 0
$ws.Range("N67").Value = 0
$ws.Range("O67").Value = 0
$ws.Range("P67").Value = 0
$ws.Range("Q67").Value = 0
$ws.Range("R67").Value = ''

# Row 68: A 33950-2023 (was row 120)
$ws.Range("A68").Value = 'A 33950-2023'
$ws.Range("B68").Value = 45134
$ws.Range("C68").Value = 46077
$ws.Range("D68").Value = 'VÄSTRA GÖTALANDS LÄN'
$ws.Range("E68").Value = 'STENUNGSUND'
$ws.Range("G68").Value = 1.7
$ws.Range("H68").Value = 0
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("K68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("M68").Value = 0
$ws.Range("N68").Value = 0
$ws.Range("O68").Value = 0
$ws.Range("P68").Value = 0
$ws.Range("Q68").Value = 0
$ws.Range("R68").Value = ''

# Row 69: A 38994-2025 (was row 52)
$ws.Range("A69").Value = 'A 38994-2025'
$ws.Range("B69").Value = 45888.33225694444
$ws.Range("C69").Value = 46077
$ws.Range("D69").Value = 'VÄSTRA GÖTALANDS LÄN'
$ws.Range("E69").Value = 'STENUNGSUND'
$ws.Range("G69").Value = 0.4
$ws.Range("H69").Value = 0
$ws.Range("I69").Value = 0
$ws.Range("J69").Value = 0
$ws.Range("K69").Value = 0
$ws.Range("L69").Value = 0
$ws.Range("M69").Value = 0
$ws.Range("N69").Value = 0
$ws.Range("O69").Value = 0
$ws.Range("P69").Value = 0
$ws.Range("Q69").Value = 0
$ws.Range("R69").Value = ''

# Row 70: A 40950-2023 (was row 61)
$ws.Range("A70").Value = 'A 40950-2023'
$ws.Range("B70").Value = 45173
$ws.Range("C70").Value = 46077
$ws.Range("D70").Value = 'VÄSTRA GÖTALANDS LÄN'
$ws.Range("E70").Value = 'STENUNGSUND'
$ws.Range("G70").Value = 6.3
$ws.Range("H70").Value = 0
$ws.Range("I70").Value = 0
$ws.Range("J70").Value = 0
$ws.Range("K70").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("M70").Value = 0
$ws.Range("N70").Value = 0
$ws.Range("O70").Value = 0
$ws.Range("P70").Value = 0
$ws.Range("Q70").Value = 0
$ws.Range("R70").Value = ''

# Row 71: A 51661-2023 (was row 107)
$ws.Range("A71").Value = 'A 51661-2023'
$ws.Range("B71").Value = 45215
$ws.Range("C71").Value = 46077
$ws.Range("D71").Value = 'VÄSTRA GÖTALANDS LÄN'
$ws.Range("E71").Value = 'STENUNGSUND'
$ws.Range("G71").Value = 9.1
$ws.Range("H71").Value = 0
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("K71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("M71").Value = 0
$ws.Range("N71").Value = 0
$ws.Range("O71").Value = 0
$ws.Range("P71").Value = 0
$ws.Range("Q71").Value = 0
$ws.Range("R71").Value = ''

# Row 72: A 19296-2025 (was row 119)
$ws.Range("A72").Value = 'A 19296-2025'
$ws.Range("B72").Value = 45769.5621875
$ws.Range("C72").Value = 46077
$ws.Range("D72").Value = 'VÄSTRA GÖTALANDS LÄN'
$ws.Range("E72").Value = 'STENUNGSUND'
$ws.Range("G72").Value = 3.3
$ws.Range("H72").Value = 0
$ws.Range("I72").Value = 0
$ws.Range("J72").Value = 0
$ws.Range("K72").Value = 0
$ws.Range("L72").Value = 0
$ws.Range("M72").Value = 0
$ws.Range("N72").Value = 0
$ws.Range("O72").Value = 0
$ws.Range("P72").Value = 0
$ws.Range("Q72").Value = 0
$ws.Range("R72").Value = ''

# Row 73: A 8948-2022 (was row 56)
$ws.Range("A73").Value = 'A 8948-2022'
$ws.Range("B73").Value = 44614
$ws.Range("C73").Value = 46077
$ws.Range("D73").Value = 'VÄSTRA GÖTALANDS LÄN'
$ws.Range("E73").Value = 'STENUNGSUND'
$ws.Range("G73").Value = 0.6
$ws.Range("H73").Value = 0
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 0
$ws.Range("K73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("M73").Value = 0
$ws.Range("N73").Value = 0
$ws.Range("O73").Value = 0
$ws.Range("P73").Value = 0
$ws.Range("Q73").Value = 0
$ws.Range("R73").Value = ''

# Row 74: A 6147-2022 (was row 94)
$ws.Range("A74").Value = 'A 6147-2022'
$ws.Range("B74").Value = 44599
$ws.Range("C74").Value = 46077
$ws.Range("D74").Value = 'VÄSTRA GÖTALANDS LÄN'
$ws.Range("E74").Value = 'STENUNGSUND'
$ws.Range("G74").Value = 2.2
$ws.Range("H74").Value = 0
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 0
$ws.Range("L74").Value = 0
$ws.Range("M74").Value = 0
$ws.Range("N74").Value = 0
$ws.Range("O74").Value = 0
$ws.Range("P74").Value = 0
$ws.Range("Q74").Value = 0
$ws.Range("R74").Value = ''

# Row 75: A 42791-2025 (was row 62)
$ws.Range("A75").Value = 'A 42791-2025'
$ws.Range("B75").Value = 45908.52364583333
$ws.Range("C75").Value = 46077
$ws.Range("D75").Value = 'VÄSTRA GÖTALANDS LÄN'
$ws.Range("E75").Value = 'STENUNGSUND'
$ws.Range("G75").Value = 1.9
$ws.Range("H75").Value = 0
$ws.Range("I75").Value = 0
$ws.Range("J75").Value = 0
$ws.Range("K75").Value = 0
$ws.Range("L75").Value = 0
$ws.Range("M75").Value = 0
$ws.Range("N75").Value = 0
$ws.Range("O75").Value = 0
$ws.Range("P75").Value = 0
$ws.Range("Q75").Value = 0
$ws.Range("R75").Value = ''

# Row 76: A 24807-2023 (was row 43)
$ws.Range("A76").Value = 'A 24807-2023'
$ws.Range("B76").Value = 45084
$ws.Range("C76").Value = 46077
$ws.Range("D76").Value = 'VÄSTRA GÖTALANDS LÄN'
$ws.Range("E76").Value = 'STENUNGSUND'
$ws.Range("G76").Value = 3
$ws.Range("H76").Value = 0
$ws.Range("I76").Value = 0
$ws.Range("J76").Value = 0
$ws.Range("K76").Value = 0
$ws.Range("L76").Value = 0
$ws.Range("M76").Value = 0
$ws.Range("N76").Value = 0
$ws.Range("O76").Value = 0
$ws.Range("P76").Value = 0
$ws.Range("Q76").Value = 0
$ws.Range("R76").Value = ''

# Row 77: A 34253-2023 (was row 89)
$ws.Range("A77").Value = 'A 34253-2023'
$ws.Range("B77").Value = 45138.5799537037
$ws.Range("C77").Value = 46077
$ws.Range("D77").Value = 'VÄSTRA GÖTALANDS LÄN'
$ws.Range("E77").Value = 'STENUNGSUND'
$ws.Range("G77").Value = 1.8
$ws.Range("H77").Value = 0
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 0
$ws.Range("L77").Value = 0
$ws.Range("M77").Value = 0
$ws.Range("N77").Value = 0
$ws.Range("O77").Value = 0
$ws.Range("P77").Value = 0
$ws.Range("Q77").Value = 0
$ws.Range("R77").Value = ''

# Row 78: A 18591-2021 (was row 45)
$ws.Range("A78").Value = 'A 18591-2021'
$ws.Range("B78").Value = 44306
$ws.Range("C78").Value = 46077
$ws.Range("D78").Value = 'VÄSTRA GÖTALANDS LÄN'
$ws.Range("E78").Value = 'STENUNGSUND'
$ws.Range("G78").Value = 2.7
$ws.Range("H78").Value = 0
$ws.Range("I78").Value = 0
$ws.Range("J78").Value = 0
$ws.Range("K78").Value = 0
$ws.Range("L78").Value = 0
$ws.Range("M78").Value = 0
$ws.Range("N78").Value = 0
$ws.Range("O78").Value = 0
$ws.Range("P78").Value = 0
$ws.Range("Q78").Value = 0
$ws.Range("R78").Value = ''

# Row 79: A 51367-2024 (was row 104)
$ws.Range("A79").Value = 'A 51367-2024'
$ws.Range("B79").Value = 45604.37143518519
$ws.Range("C79").Value = 46077
$ws.Range("D79").Value = 'VÄSTRA GÖTALANDS LÄN'
$ws.Range("E79").Value = 'STENUNGSUND'
$ws.Range("G79").Value = 3
$ws.Range("H79").Value = 0
$ws.Range("I79").Value = 0
$ws.Range("J79").Value = 0
$ws.Range("K79").Value = 0
$ws.Range("L79").Value = 0
$ws.Range("M79").Value = 0
$ws.Range("N79").Value = 0
$ws.Range("O79").Value = 0
$ws.Range("P79").Value = 0
$ws.Range("Q79").Value = 0
$ws.Range("R79").Value = ''

# Row 80: A 13586-2024 (was row 42)
$ws.Range("A80").Value = 'A 13586-2024'
$ws.Range("B80").Value = 45389.85628472222
$ws.Range("C80").Value = 46077
$ws.Range("D80").Value = 'VÄSTRA GÖTALANDS LÄN'
$ws.Range("E80").Value = 'STENUNGSUND'
$ws.Range("G80").Value = 0.8
$ws.Range("H80").Value = 0
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 0
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 0
$ws.Range("M80").Value = 0
$ws.Range("N80").Value = 0
$ws.Range("O80").Value = 0
$ws.Range("P80").Value = 0
$ws.Range("Q80").Value = 0
$ws.Range("R80").Value = ''

# Row 81: A 54931-2024 (was row 98)
$ws.Range("A81").Value = 'A 54931-2024'
$ws.Range("B81").Value = 45618.74254629629
$ws.Range("C81").Value = 46077
$ws.Range("D81").Value = 'VÄSTRA GÖTALANDS LÄN'
$ws.Range("E81").Value = 'STENUNGSUND'
$ws.Range("G81").Value = 0.5
$ws.Range("H81").Value = 0
$ws.Range("I81").Value = 0
$ws.Range("J81").Value = 0
$ws.Range("K81").Value = 0
$ws.Range("L81").Value = 0
$ws.Range("M81").Value = 0
$ws.Range("N81").Value = 0
$ws.Range("O81").Value = 0
$ws.Range("P81").Value = 0
$ws.Range("Q81").Value = 0
$ws.Range("R81").Value = ''

# Row 82: A 45710-2025 (was row 64)
$ws.Range("A82").Value = 'A 45710-2025'
$ws.Range("B82").Value = 45923.4296875
$ws.Range("C82").Value = 46077
$ws.Range("D82").Value = 'VÄSTRA GÖTALANDS LÄN'
$ws.Range("E82").Value = 'STENUNGSUND'
$ws.Range("G82").Value = 1.3
$ws.Range("H82").Value = 0
$ws.Range("I82").Value = 0
$ws.Range("J82").Value = 0
$ws.Range("K82").Value = 0
$ws.Range("L82").Value = 0
$ws.Range("M82").Value = 0
$ws.Range("N82").Value = 0
$ws.Range("O82").Value = 0
$ws.Range("P82").Value = 0
$ws.Range("Q82").Value = 0
$ws.Range("R82").Value = ''

# Row 83: A 22728-2025 (was row 73)
$ws.Range("A83").Value = 'A 22728-2025'
$ws.Range("B83").Value = 45789.57516203704
$ws.Range("C83").Value = 46077
$ws.Range("D83").Value = 'VÄSTRA GÖTALANDS LÄN'
$ws.Range("E83").Value = 'STENUNGSUND'
$ws.Range("G83").Value = 1.3
$ws.Range("H83").Value = 0
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 0
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 0
$ws.Range("M83").Value = 0
$ws.Range("N83").Value = 0
$ws.Range("O83").Value = 0
$ws.Range("P83").Value = 0
$ws.Range("Q83").Value = 0
$ws.Range("R83").Value = ''

# Row 84: A 45000-2024 (was row 34)
$ws.Range("A84").Value = 'A 45000-2024'
$ws.Range("B84").Value = 45575
$ws.Range("C84").Value = 46077
$ws.Range("D84").Value = 'VÄSTRA GÖTALANDS LÄN'
$ws.Range("E84").Value = 'STENUNGSUND'
$ws.Range("G84").Value = 1.2
$ws.Range("H84").Value = 0
$ws.Range("I84").Value = 0
$ws.Range("J84").Value = 0
$ws.Range("K84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("M84").Value = 0
$ws.Range("N84").Value = 0
$ws.Range("O84").Value = 0
$ws.Range("P84").Value = 0
$ws.Range("Q84").Value = 0
$ws.Range("R84").Value = ''

# Row 85: A 57766-2024 (was row 67)
$ws.Range("A85").Value = 'A 57766-2024'
$ws.Range("B85").Value = 45630.87388888889
$ws.Range("C85").Value = 46077
$ws.Range("D85").Value = 'VÄSTRA GÖTALANDS LÄN'
$ws.Range("E85").Value = 'STENUNGSUND'
$ws.Range("G85").Value = 2.5
$ws.Range("H85").Value = 0
$ws.Range("I85").Value = 0
$ws.Range("J85").Value = 0
$ws.Range("K85").Value = 0
$ws.Range("L85").Value = 0
$ws.Range("M85").Value = 0
$ws.Range("N85").Value = 0
$ws.Range("O85").Value = 0
$ws.Range("P85").Value = 0
$ws.Range("Q85").Value = 0
$ws.Range("R85").Value = ''

# Row 86: A 57767-2024 (was row 68)
$ws.Range("A86").Value = 'A 57767-2024'
$ws.Range("B86").Value = 45630.87564814815
$ws.Range("C86").Value = 46077
$ws.Range("D86").Value = 'VÄSTRA GÖTALANDS LÄN'
$ws.Range("E86").Value = 'STENUNGSUND'
$ws.Range("G86").Value = 1.3
$ws.Range("H86").Value = 0
$ws.Range("I86").Value = 0
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 0
$ws.Range("L86").Value = 0
$ws.Range("M86").Value = 0
$ws.Range("N86").Value = 0
$ws.Range("O86").Value = 0
$ws.Range("P86").Value = 0
$ws.Range("Q86").Value = 0
$ws.Range("R86").Value = ''

# Row 87: A 57769-2024 (was row 69)
$ws.Range("A87").Value = 'A 57769-2024'
$ws.Range("B87").Value = 45630.88071759259
$ws.Range("C87").Value = 46077
$ws.Range("D87").Value = 'VÄSTRA GÖTALANDS LÄN'
$ws.Range("E87").Value = 'STENUNGSUND'
$ws.Range("G87").Value = 0.6
$ws.Range("H87").Value = 0
$ws.Range("I87").Value = 0
$ws.Range("J87").Value = 0
$ws.Range("K87").Value = 0
$ws.Range("L87").Value = 0
$ws.Range("M87").Value = 0
$ws.Range("N87").Value = 0
$ws.Range("O87").Value = 0
$ws.Range("P87").Value = 0
$ws.Range("Q87").Value = 0
$ws.Range("R87").Value = ''

# Row 88: A 57768-2024 (was row 66)
$ws.Range("A88").Value = 'A 57768-2024'
$ws.Range("B88").Value = 45630.87756944444
$ws.Range("C88").Value = 46077
$ws.Range("D88").Value = 'VÄSTRA GÖTALANDS LÄN'
$ws.Range("E88").Value = 'STENUNGSUND'
$ws.Range("G88").Value = 1
$ws.Range("H88").Value = 0
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 0
$ws.Range("K88").Value = 0
$ws.Range("L88").Value = 0
$ws.Range("M88").Value = 0
$ws.Range("N88").Value = 0
$ws.Range("O88").Value = 0
$ws.Range("P88").Value = 0
$ws.Range("Q88").Value = 0
$ws.Range("R88").Value = ''

# Row 89: A 3569-2023 (was row 55)
$ws.Range("A89").Value = 'A 3569-2023'
$ws.Range("B89").Value = 44950.48391203704
$ws.Range("C89").Value = 46077
$ws.Range("D89").Value = 'VÄSTRA GÖTALANDS LÄN'
$ws.Range("E89").Value = 'STENUNGSUND'
$ws.Range("G89").Value = 2.2
$ws.Range("H89").Value = 0
$ws.Range("I89").Value = 0
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 0
$ws.Range("L89").Value = 0
$ws.Range("M89").Value = 0
$ws.Range("N89").Value = 0
$ws.Range("O89").Value = 0
$ws.Range("P89").Value = 0
$ws.Range("Q89").Value = 0
$ws.Range("R89").Value = ''

# Row 90: A 8684-2024 (was row 36)
$ws.Range("A90").Value = 'A 8684-2024'
$ws.Range("B90").Value = 45356.36320601852
$ws.Range("C90").Value = 46077
$ws.Range("D90").Value = 'VÄSTRA GÖTALANDS LÄN'
$ws.Range("E90").Value = 'STENUNGSUND'
$ws.Range("G90").Value = 2.3
$ws.Range("H90").Value = 0
$ws.Range("I90").Value = 0
$ws.Range("J90").Value = 0
$ws.Range("K90").Value = 0
$ws.Range("L90").Value = 0
$ws.Range("M90").Value = 0
$ws.Range("N90").Value = 0
$ws.Range("O90").Value = 0
$ws.Range("P90").Value = 0
$ws.Range("Q90").Value = 0
$ws.Range("R90").Value = ''

# Row 91: A 8544-2023 (was row 57)
$ws.Range("A91").Value = 'A 8544-2023'
$ws.Range("B91").Value = 44977
$ws.Range("C91").Value = 46077
$ws.Range("D91").Value = 'VÄSTRA GÖTALANDS LÄN'
$ws.Range("E91").Value = 'STENUNGSUND'
$ws.Range("G91").Value = 8.199999999999999
$ws.Range("H91").Value = 0
$ws.Range("I91").Value = 0
$ws.Range("J91").Value = 0
$ws.Range("K91").Value = 0
$ws.Range("L91").Value = 0
$ws.Range("M91").Value = 0
$ws.Range("N91").Value = 0
$ws.Range("O91").Value = 0
$ws.Range("P91").Value = 0
$ws.Range("Q91").Value = 0
$ws.Range("R91").Value = ''

# Row 92: A 60550-2023 (was row 74)
$ws.Range("A92").Value = 'A 60550-2023'
$ws.Range("B92").Value = 45259.84060185185
$ws.Range("C92").Value = 46077
$ws.Range("D92").Value = 'VÄSTRA GÖTALANDS LÄN'
$ws.Range("E92").Value = 'STENUNGSUND'
$ws.Range("G92").Value = 0.5
$ws.Range("H92").Value = 0
$ws.Range("I92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("K92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("M92").Value = 0
$ws.Range("N92").Value = 0
$ws.Range("O92").Value = 0
$ws.Range("P92").Value = 0
$ws.Range("Q92").Value = 0
$ws.Range("R92").Value = ''

# Row 93: A 49621-2025 (was row 77)
$ws.Range("A93").Value = 'A 49621-2025'
$ws.Range("B93").Value = 45939.59266203704
$ws.Range("C93").Value = 46077
$ws.Range("D93").Value = 'VÄSTRA GÖTALANDS LÄN'
$ws.Range("E93").Value = 'STENUNGSUND'
$ws.Range("G93").Value = 3.3
$ws.Range("H93").Value = 0
$ws.Range("I93").Value = 0
$ws.Range("J93").Value = 0
$ws.Range("K93").Value = 0
$ws.Range("L93").Value = 0
$ws.Range("M93").Value = 0
$ws.Range("N93").Value = 0
$ws.Range("O93").Value = 0
$ws.Range("P93").Value = 0
$ws.Range("Q93").Value = 0
$ws.Range("R93").Value = ''

# Row 94: A 50205-2025 (was row 78)
$ws.Range("A94").Value = 'A 50205-2025'
$ws.Range("B94").Value = 45943.65840277778
$ws.Range("C94").Value = 46077
$ws.Range("D94").Value = 'VÄSTRA GÖTALANDS LÄN'
$ws.Range("E94").Value = 'STENUNGSUND'
$ws.Range("G94").Value = 3.6
$ws.Range("H94").Value = 0
$ws.Range("I94").Value = 0
$ws.Range("J94").Value = 0
$ws.Range("K94").Value = 0
$ws.Range("L94").Value = 0
$ws.Range("M94").Value = 0
$ws.Range("N94").Value = 0
$ws.Range("O94").Value = 0
$ws.Range("P94").Value = 0
$ws.Range("Q94").Value = 0
$ws.Range("R94").Value = ''

# Row 95: A 50491-2025 (was row 81)
$ws.Range("A95").Value = 'A 50491-2025'
$ws.Range("B95").Value = 45945
$ws.Range("C95").Value = 46077
$ws.Range("D95").Value = 'VÄSTRA GÖTALANDS LÄN'
$ws.Range("E95").Value = 'STENUNGSUND'
$ws.Range("G95").Value = 0.9
$ws.Range("H95").Value = 0
$ws.Range("I95").Value = 0
$ws.Range("J95").Value = 0
$ws.Range("K95").Value = 0
$ws.Range("L95").Value = 0
$ws.Range("M95").Value = 0
$ws.Range("N95").Value = 0
$ws.Range("O95").Value = 0
$ws.Range("P95").Value = 0
$ws.Range("Q95").Value = 0
$ws.Range("R95").Value = ''

# Row 96: A 50508-2025 (was row 82)
$ws.Range("A96").Value = 'A 50508-2025'
$ws.Range("B96").Value = 45945.43429398148
$ws.Range("C96").Value = 46077
$ws.Range("D96").Value = 'VÄSTRA GÖTALANDS LÄN'
$ws.Range("E96").Value = 'STENUNGSUND'
$ws.Range("G96").Value = 3.6
$ws.Range("H96").Value = 0
$ws.Range("I96").Value = 0
$ws.Range("J96").Value = 0
$ws.Range("K96").Value = 0
$ws.Range("L96").Value = 0
$ws.Range("M96").Value = 0
$ws.Range("N96").Value = 0
$ws.Range("O96").Value = 0
$ws.Range("P96").Value = 0
$ws.Range("Q96").Value = 0
$ws.Range("R96").Value = ''

# Row 97: A 50543-2025 (was row 85)
$ws.Range("A97").Value = 'A 50543-2025'
$ws.Range("B97").Value = 45945.50262731482
$ws.Range("C97").Value = 46077
$ws.Range("D97").Value = 'VÄSTRA GÖTALANDS LÄN'
$ws.Range("E97").Value = 'STENUNGSUND'
$ws.Range("G97").Value = 0.8
$ws.Range("H97").Value = 0
$ws.Range("I97").Value = 0
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 0
$ws.Range("L97").Value = 0
$ws.Range("M97").Value = 0
$ws.Range("N97").Value = 0
$ws.Range("O97").Value = 0
$ws.Range("P97").Value = 0
$ws.Range("Q97").Value = 0
$ws.Range("R97").Value = ''

# Row 98: A 50538-2025 (was row 84)
$ws.Range("A98").Value = 'A 50538-2025'
$ws.Range("B98").Value = 45945.49040509259
$ws.Range("C98").Value = 46077
$ws.Range("D98").Value = 'VÄSTRA GÖTALANDS LÄN'
$ws.Range("E98").Value = 'STENUNGSUND'
$ws.Range("G98").Value = 3.8
$ws.Range("H98").Value = 0
$ws.Range("I98").Value = 0
$ws.Range("J98").Value = 0
$ws.Range("K98").Value = 0
$ws.Range("L98").Value = 0
$ws.Range("M98").Value = 0
$ws.Range("N98").Value = 0
$ws.Range("O98").Value = 0
$ws.Range("P98").Value = 0
$ws.Range("Q98").Value = 0
$ws.Range("R98").Value = ''

# Row 99: A 51450-2025 (was row 87)
$ws.Range("A99").Value = 'A 51450-2025'
$ws.Range("B99").Value = 45950
$ws.Range("C99").Value = 46077
$ws.Range("D99").Value = 'VÄSTRA GÖTALANDS LÄN'
$ws.Range("E99").Value = 'STENUNGSUND'
$ws.Range("G99").Value = 0.9
$ws.Range("H99").Value = 0
$ws.Range("I99").Value = 0
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 0
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = 0
$ws.Range("N99").Value = 0
$ws.Range("O99").Value = 0
$ws.Range("P99").Value = 0
$ws.Range("Q99").Value = 0
$ws.Range("R99").Value = ''

# Row 100: A 51546-2025 (was row 86)
$ws.Range("A100").Value = 'A 51546-2025'
$ws.Range("B100").Value = 45950.7683912037
$ws.Range("C100").Value = 46077
$ws.Range("D100").Value = 'VÄSTRA GÖTALANDS LÄN'
$ws.Range("E100").Value = 'STENUNGSUND'
$ws.Range("G100").Value = 1.2
$ws.Range("H100").Value = 0
$ws.Range("I100").Value = 0
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 0
$ws.Range("L100").Value = 0
$ws.Range("M100").Value = 0
$ws.Range("N100").Value = 0
$ws.Range("O100").Value = 0
$ws.Range("P100").Value = 0
$ws.Range("Q100").Value = 0
$ws.Range("R100").Value = ''

# Row 101: A 52012-2025 (was row 90)
$ws.Range("A101").Value = 'A 52012-2025'
$ws.Range("B101").Value = 45952.59604166666
$ws.Range("C101").Value = 46077
$ws.Range("D101").Value = 'VÄSTRA GÖTALANDS LÄN'
$ws.Range("E101").Value = 'STENUNGSUND'
$ws.Range("G101").Value = 1.2
$ws.Range("H101").Value = 0
$ws.Range("I101").Value = 0
$ws.Range("J101").Value = 0
$ws.Range("K101").Value = 0
$ws.Range("L101").Value = 0
$ws.Range("M101").Value = 0
$ws.Range("N101").Value = 0
$ws.Range("O101").Value = 0
$ws.Range("P101").Value = 0
$ws.Range("Q101").Value = 0
$ws.Range("R101").Value = ''

# Row 102: A 51992-2025 (was row 88)
$ws.Range("A102").Value = 'A 51992-2025'
$ws.Range("B102").Value = 45952.58126157407
$ws.Range("C102").Value = 46077
$ws.Range("D102").Value = 'VÄSTRA GÖTALANDS LÄN'
$ws.Range("E102").Value = 'STENUNGSUND'
$ws.Range("G102").Value = 4.7
$ws.Range("H102").Value = 0
$ws.Range("I102").Value = 0
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 0
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = 0
$ws.Range("N102").Value = 0
$ws.Range("O102").Value = 0
$ws.Range("P102").Value = 0
$ws.Range("Q102").Value = 0
$ws.Range("R102").Value = ''

# Row 103: A 8567-2023 (was row 93)
$ws.Range("A103").Value = 'A 8567-2023'
$ws.Range("B103").Value = 44977.63353009259
$ws.Range("C103").Value = 46077
$ws.Range("D103").Value = 'VÄSTRA GÖTALANDS LÄN'
$ws.Range("E103").Value = 'STENUNGSUND'
$ws.Range("G103").Value = 0.9
$ws.Range("H103").Value = 0
$ws.Range("I103").Value = 0
$ws.Range("J103").Value = 0
$ws.Range("K103").Value = 0
$ws.Range("L103").Value = 0
$ws.Range("M103").Value = 0
$ws.Range("N103").Value = 0
$ws.Range("O103").Value = 0
$ws.Range("P103").Value = 0
$ws.Range("Q103").Value = 0
$ws.Range("R103").Value = ''

# Row 104: A 54247-2025 (was row 92)
$ws.Range("A104").Value = 'A 54247-2025'
$ws.Range("B104").Value = 45964.64162037037
$ws.Range("C104").Value = 46077
$ws.Range("D104").Value = 'VÄSTRA GÖTALANDS LÄN'
$ws.Range("E104").Value = 'STENUNGSUND'
$ws.Range("G104").Value = 1
$ws.Range("H104").Value = 0
$ws.Range("I104").Value = 0
$ws.Range("J104").Value = 0
$ws.Range("K104").Value = 0
$ws.Range("L104").Value = 0
$ws.Range("M104").Value = 0
$ws.Range("N104").Value = 0
$ws.Range("O104").Value = 0
$ws.Range("P104").Value = 0
$ws.Range("Q104").Value = 0
$ws.Range("R104").Value = ''

# Row 105: A 15724-2025 (was row 114)
$ws.Range("A105").Value = 'A 15724-2025'
$ws.Range("B105").Value = 45748
$ws.Range("C105").Value = 46077
$ws.Range("D105").Value = 'VÄSTRA GÖTALANDS LÄN'
$ws.Range("E105").Value = 'STENUNGSUND'
$ws.Range("G105").Value = 3.8
$ws.Range("H105").Value = 0
$ws.Range("I105").Value = 0
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 0
$ws.Range("L105").Value = 0
$ws.Range("M105").Value = 0
$ws.Range("N105").Value = 0
$ws.Range("O105").Value = 0
$ws.Range("P105").Value = 0
$ws.Range("Q105").Value = 0
$ws.Range("R105").Value = ''

# Row 106: A 15726-2025 (was row 115)
$ws.Range("A106").Value = 'A 15726-2025'
$ws.Range("B106").Value = 45748
$ws.Range("C106").Value = 46077
$ws.Range("D106").Value = 'VÄSTRA GÖTALANDS LÄN'
$ws.Range("E106").Value = 'STENUNGSUND'
$ws.Range("G106").Value = 1.2
$ws.Range("H106").Value = 0
$ws.Range("I106").Value = 0
$ws.Range("J106").Value = 0
$ws.Range("K106").Value = 0
$ws.Range("L106").Value = 0
$ws.Range("M106").Value = 0
$ws.Range("N106").Value = 0
$ws.Range("O106").Value = 0
$ws.Range("P106").Value = 0
$ws.Range("Q106").Value = 0
$ws.Range("R106").Value = ''

# Row 107: A 36018-2025 (was row 95)
$ws.Range("A107").Value = 'A 36018-2025'
$ws.Range("B107").Value = 45863.66916666667
$ws.Range("C107").Value = 46077
$ws.Range("D107").Value = 'VÄSTRA GÖTALANDS LÄN'
$ws.Range("E107").Value = 'STENUNGSUND'
$ws.Range("G107").Value = 1.3
$ws.Range("H107").Value = 0
$ws.Range("I107").Value = 0
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 0
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = 0
$ws.Range("N107").Value = 0
$ws.Range("O107").Value = 0
$ws.Range("P107").Value = 0
$ws.Range("Q107").Value = 0
$ws.Range("R107").Value = ''

# Row 108: A 809-2025 (was row 79)
$ws.Range("A108").Value = 'A 809-2025'
$ws.Range("B108").Value = 45665.5518287037
$ws.Range("C108").Value = 46077
$ws.Range("D108").Value = 'VÄSTRA GÖTALANDS LÄN'
$ws.Range("E108").Value = 'STENUNGSUND'
$ws.Range("G108").Value = 1
$ws.Range("H108").Value = 0
$ws.Range("I108").Value = 0
$ws.Range("J108").Value = 0
$ws.Range("K108").Value = 0
$ws.Range("L108").Value = 0
$ws.Range("M108").Value = 0
$ws.Range("N108").Value = 0
$ws.Range("O108").Value = 0
$ws.Range("P108").Value = 0
$ws.Range("Q108").Value = 0
$ws.Range("R108").Value = ''

# Row 109: A 22187-2025 (was row 121)
$ws.Range("A109").Value = 'A 22187-2025'
$ws.Range("B109").Value = 45785.62372685185
$ws.Range("C109").Value = 46077
$ws.Range("D109").Value = 'VÄSTRA GÖTALANDS LÄN'
$ws.Range("E109").Value = 'STENUNGSUND'
$ws.Range("G109").Value = 2.1
$ws.Range("H109").Value = 0
$ws.Range("I109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("K109").Value = 0
$ws.Range("L109").Value = 0
$ws.Range("M109").Value = 0
$ws.Range("N109").Value = 0
$ws.Range("O109").Value = 0
$ws.Range("P109").Value = 0
$ws.Range("Q109").Value = 0
$ws.Range("R109").Value = ''

# Row 110: A 60687-2025 (was row 103)
$ws.Range("A110").Value = 'A 60687-2025'
$ws.Range("B110").Value = 45996.5859837963
$ws.Range("C110").Value = 46077
$ws.Range("D110").Value = 'VÄSTRA GÖTALANDS LÄN'
$ws.Range("E110").Value = 'STENUNGSUND'
$ws.Range("G110").Value = 3.9
$ws.Range("H110").Value = 0
$ws.Range("I110").Value = 0
$ws.Range("J110").Value = 0
$ws.Range("K110").Value = 0
$ws.Range("L110").Value = 0
$ws.Range("M110").Value = 0
$ws.Range("N110").Value = 0
$ws.Range("O110").Value = 0
$ws.Range("P110").Value = 0
$ws.Range("Q110").Value = 0
$ws.Range("R110").Value = ''

# Row 111: A 2472-2026 (was row 101)
$ws.Range("A111").Value = 'A 2472-2026'
$ws.Range("B111").Value = 46037.37466435185
$ws.Range("C111").Value = 46077
$ws.Range("D111").Value = 'VÄSTRA GÖTALANDS LÄN'
$ws.Range("E111").Value = 'STENUNGSUND'
$ws.Range("G111").Value = 1.1
$ws.Range("H111").Value = 0
$ws.Range("I111").Value = 0
$ws.Range("J111").Value = 0
$ws.Range("K111").Value = 0
$ws.Range("L111").Value = 0
$ws.Range("M111").Value = 0
$ws.Range("N111").Value = 0
$ws.Range("O111").Value = 0
$ws.Range("P111").Value = 0
$ws.Range("Q111").Value = 0
$ws.Range("R111").Value = ''

# Row 112: A 13116-2025 (was row 123)
$ws.Range("A112").Value = 'A 13116-2025'
$ws.Range("B112").Value = 45734.71525462963
$ws.Range("C112").Value = 46077
$ws.Range("D112").Value = 'VÄSTRA GÖTALANDS LÄN'
$ws.Range("E112").Value = 'STENUNGSUND'
$ws.Range("G112").Value = 1.5
$ws.Range("H112").Value = 0
$ws.Range("I112").Value = 0
$ws.Range("J112").Value = 0
$ws.Range("K112").Value = 0
$ws.Range("L112").Value = 0
$ws.Range("M112").Value = 0
$ws.Range("N112").Value = 0
$ws.Range("O112").Value = 0
$ws.Range("P112").Value = 0
$ws.Range("Q112").Value = 0
$ws.Range("R112").Value = ''

# Row 113: A 2481-2026 (was row 100)
$ws.Range("A113").Value = 'A 2481-2026'
$ws.Range("B113").Value = 46037.38494212963
$ws.Range("C113").Value = 46077
$ws.Range("D113").Value = 'VÄSTRA GÖTALANDS LÄN'
$ws.Range("E113").Value = 'STENUNGSUND'
$ws.Range("G113").Value = 1.4
$ws.Range("H113").Value = 0
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = 0
$ws.Range("N113").Value = 0
$ws.Range("O113").Value = 0
$ws.Range("P113").Value = 0
$ws.Range("Q113").Value = 0
$ws.Range("R113").Value = ''

# Row 114: A 9080-2025 (was row 122)
$ws.Range("A114").Value = 'A 9080-2025'
$ws.Range("B114").Value = 45713.66313657408
$ws.Range("C114").Value = 46077
$ws.Range("D114").Value = 'VÄSTRA GÖTALANDS LÄN'
$ws.Range("E114").Value = 'STENUNGSUND'
$ws.Range("G114").Value = 1.8
$ws.Range("H114").Value = 0
$ws.Range("I114").Value = 0
$ws.Range("J114").Value = 0
$ws.Range("K114").Value = 0
$ws.Range("L114").Value = 0
$ws.Range("M114").Value = 0
$ws.Range("N114").Value = 0
$ws.Range("O114").Value = 0
$ws.Range("P114").Value = 0
$ws.Range("Q114").Value = 0
$ws.Range("R114").Value = ''

# Row 115: A 12589-2025 (was row 124)
$ws.Range("A115").Value = 'A 12589-2025'
$ws.Range("B115").Value = 45732.7803125
$ws.Range("C115").Value = 46077
$ws.Range("D115").Value = 'VÄSTRA GÖTALANDS LÄN'
$ws.Range("E115").Value = 'STENUNGSUND'
$ws.Range("G115").Value = 2
$ws.Range("H115").Value = 0
$ws.Range("I115").Value = 0
$ws.Range("J115").Value = 0
$ws.Range("K115").Value = 0
$ws.Range("L115").Value = 0
$ws.Range("M115").Value = 0
$ws.Range("N115").Value = 0
$ws.Range("O115").Value = 0
$ws.Range("P115").Value = 0
$ws.Range("Q115").Value = 0
$ws.Range("R115").Value = ''

# Row 116: A 14906-2025 (was row 125)
$ws.Range("A116").Value = 'A 14906-2025'
$ws.Range("B116").Value = 45743.49373842592
$ws.Range("C116").Value = 46077
$ws.Range("D116").Value = 'VÄSTRA GÖTALANDS LÄN'
$ws.Range("E116").Value = 'STENUNGSUND'
$ws.Range("G116").Value = 1.1
$ws.Range("H116").Value = 0
$ws.Range("I116").Value = 0
$ws.Range("J116").Value = 0
$ws.Range("K116").Value = 0
$ws.Range("L116").Value = 0
$ws.Range("M116").Value = 0
$ws.Range("N116").Value = 0
$ws.Range("O116").Value = 0
$ws.Range("P116").Value = 0
$ws.Range("Q116").Value = 0
$ws.Range("R116").Value = ''

# Row 117: A 8961-2025 (was row 126)
$ws.Range("A117").Value = 'A 8961-2025'
$ws.Range("B117").Value = 45713.48292824074
$ws.Range("C117").Value = 46077
$ws.Range("D117").Value = 'VÄSTRA GÖTALANDS LÄN'
$ws.Range("E117").Value = 'STENUNGSUND'
$ws.Range("G117").Value = 2
$ws.Range("H117").Value = 0
$ws.Range("I117").Value = 0
$ws.Range("J117").Value = 0
$ws.Range("K117").Value = 0
$ws.Range("L117").Value = 0
$ws.Range("M117").Value = 0
$ws.Range("N117").Value = 0
$ws.Range("O117").Value = 0
$ws.Range("P117").Value = 0
$ws.Range("Q117").Value = 0
$ws.Range("R117").Value = ''

# Row 118: A 11362-2025 (was row 102)
$ws.Range("A118").Value = 'A 11362-2025'
$ws.Range("B118").Value = 45726
$ws.Range("C118").Value = 46077
$ws.Range("D118").Value = 'VÄSTRA GÖTALANDS LÄN'
$ws.Range("E118").Value = 'STENUNGSUND'
$ws.Range("G118").Value = 6.4
$ws.Range("H118").Value = 0
$ws.Range("I118").Value = 0
$ws.Range("J118").Value = 0
$ws.Range("K118").Value = 0
$ws.Range("L118").Value = 0
$ws.Range("M118").Value = 0
$ws.Range("N118").Value = 0
$ws.Range("O118").Value = 0
$ws.Range("P118").Value = 0
$ws.Range("Q118").Value = 0
$ws.Range("R118").Value = ''

# Row 119: A 24211-2025 (was row 127)
$ws.Range("A119").Value = 'A 24211-2025'
$ws.Range("B119").Value = 45797.34690972222
$ws.Range("C119").Value = 46077
$ws.Range("D119").Value = 'VÄSTRA GÖTALANDS LÄN'
$ws.Range("E119").Value = 'STENUNGSUND'
$ws.Range("G119").Value = 1.4
$ws.Range("H119").Value = 0
$ws.Range("I119").Value = 0
$ws.Range("J119").Value = 0
$ws.Range("K119").Value = 0
$ws.Range("L119").Value = 0
$ws.Range("M119").Value = 0
$ws.Range("N119").Value = 0
$ws.Range("O119").Value = 0
$ws.Range("P119").Value = 0
$ws.Range("Q119").Value = 0
$ws.Range("R119").Value = ''

# Row 120: A 4769-2026 (was row 106)
$ws.Range("A120").Value = 'A 4769-2026'
$ws.Range("B120").Value = 46048.54363425926
$ws.Range("C120").Value = 46077
$ws.Range("D120").Value = 'VÄSTRA GÖTALANDS LÄN'
$ws.Range("E120").Value = 'STENUNGSUND'
$ws.Range("G120").Value = 1
$ws.Range("H120").Value = 0
$ws.Range("I120").Value = 0
$ws.Range("J120").Value = 0
$ws.Range("K120").Value = 0
$ws.Range("L120").Value = 0
$ws.Range("M120").Value = 0
$ws.Range("N120").Value = 0
$ws.Range("O120").Value = 0
$ws.Range("P120").Value = 0
$ws.Range("Q120").Value = 0
$ws.Range("R120").Value = ''

# Row 121: A 60832-2024 (was row 128)
$ws.Range("A121").Value = 'A 60832-2024'
$ws.Range("B121").Value = 45644.64081018518
$ws.Range("C121").Value = 46077
$ws.Range("D121").Value = 'VÄSTRA GÖTALANDS LÄN'
$ws.Range("E121").Value = 'STENUNGSUND'
$ws.Range("G121").Value = 2.8
$ws.Range("H121").Value = 0
$ws.Range("I121").Value = 0
$ws.Range("J121").Value = 0
$ws.Range("K121").Value = 0
$ws.Range("L121").Value = 0
$ws.Range("M121").Value = 0
$ws.Range("N121").Value = 0
$ws.Range("O121").Value = 0
$ws.Range("P121").Value = 0
$ws.Range("Q121").Value = 0
$ws.Range("R121").Value = ''

# Row 122: A 48365-2024 (was row 129)
$ws.Range("A122").Value = 'A 48365-2024'
$ws.Range("B122").Value = 45590.57652777778
$ws.Range("C122").Value = 46077
$ws.Range("D122").Value = 'VÄSTRA GÖTALANDS LÄN'
$ws.Range("E122").Value = 'STENUNGSUND'
$ws.Range("G122").Value = 2.3
$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = 0
$ws.Range("N122").Value = 0
$ws.Range("O122").Value = 0
$ws.Range("P122").Value = 0
$ws.Range("Q122").Value = 0
$ws.Range("R122").Value = ''

# Row 123: A 9893-2025 (was row 47)
$ws.Range("A123").Value = 'A 9893-2025'
$ws.Range("B123").Value = 45716.70841435185
$ws.Range("C123").Value = 46077
$ws.Range("D123").Value = 'VÄSTRA GÖTALANDS LÄN'
$ws.Range("E123").Value = 'STENUNGSUND'
$ws.Range("G123").Value = 3.2
$ws.Range("H123").Value = 0
$ws.Range("I123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("K123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("M123").Value = 0
$ws.Range("N123").Value = 0
$ws.Range("O123").Value = 0
$ws.Range("P123").Value = 0
$ws.Range("Q123").Value = 0
$ws.Range("R123").Value = ''

# Row 124: A 2412-2024 (was row 130)
$ws.Range("A124").Value = 'A 2412-2024'
$ws.Range("B124").Value = 45310
$ws.Range("C124").Value = 46077
$ws.Range("D124").Value = 'VÄSTRA GÖTALANDS LÄN'
$ws.Range("E124").Value = 'STENUNGSUND'
$ws.Range("G124").Value = 1.1
$ws.Range("H124").Value = 0
$ws.Range("I124").Value = 0
$ws.Range("J124").Value = 0
$ws.Range("K124").Value = 0
$ws.Range("L124").Value = 0
$ws.Range("M124").Value = 0
$ws.Range("N124").Value = 0
$ws.Range("O124").Value = 0
$ws.Range("P124").Value = 0
$ws.Range("Q124").Value = 0
$ws.Range("R124").Value = ''

# Row 125: A 7562-2026 (was row 110)
$ws.Range("A125").Value = 'A 7562-2026'
$ws.Range("B125").Value = 46060.45351851852
$ws.Range("C125").Value = 46077
$ws.Range("D125").Value = 'VÄSTRA GÖTALANDS LÄN'
$ws.Range("E125").Value = 'STENUNGSUND'
$ws.Range("G125").Value = 2.3
$ws.Range("H125").Value = 0
$ws.Range("I125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("K125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("M125").Value = 0
$ws.Range("N125").Value = 0
$ws.Range("O125").Value = 0
$ws.Range("P125").Value = 0
$ws.Range("Q125").Value = 0
$ws.Range("R125").Value = ''

# Row 126: A 31625-2025 (was row 60)
$ws.Range("A126").Value = 'A 31625-2025'
$ws.Range("B126").Value = 45834.29967592593
$ws.Range("C126").Value = 46077
$ws.Range("D126").Value = 'VÄSTRA GÖTALANDS LÄN'
$ws.Range("E126").Value = 'STENUNGSUND'
$ws.Range("G126").Value = 2.5
$ws.Range("H126").Value = 0
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = 0
$ws.Range("N126").Value = 0
$ws.Range("O126").Value = 0
$ws.Range("P126").Value = 0
$ws.Range("Q126").Value = 0
$ws.Range("R126").Value = ''

# Row 127: A 32997-2025 (was row 111)
$ws.Range("A127").Value = 'A 32997-2025'
$ws.Range("B127").Value = 45840.33818287037
$ws.Range("C127").Value = 46077
$ws.Range("D127").Value = 'VÄSTRA GÖTALANDS LÄN'
$ws.Range("E127").Value = 'STENUNGSUND'
$ws.Range("G127").Value = 1.6
$ws.Range("H127").Value = 0
$ws.Range("I127").Value = 0
$ws.Range("J127").Value = 0
$ws.Range("K127").Value = 0
$ws.Range("L127").Value = 0
$ws.Range("M127").Value = 0
$ws.Range("N127").Value = 0
$ws.Range("O127").Value = 0
$ws.Range("P127").Value = 0
$ws.Range("Q127").Value = 0
$ws.Range("R127").Value = ''

# Row 128: A 35172-2025 (was row 63)
$ws.Range("A128").Value = 'A 35172-2025'
$ws.Range("B128").Value = 45852.81349537037
$ws.Range("C128").Value = 46077
$ws.Range("D128").Value = 'VÄSTRA GÖTALANDS LÄN'
$ws.Range("E128").Value = 'STENUNGSUND'
$ws.Range("G128").Value = 4.5
$ws.Range("H128").Value = 0
$ws.Range("I128").Value = 0
$ws.Range("J128").Value = 0
$ws.Range("K128").Value = 0
$ws.Range("L128").Value = 0
$ws.Range("M128").Value = 0
$ws.Range("N128").Value = 0
$ws.Range("O128").Value = 0
$ws.Range("P128").Value = 0
$ws.Range("Q128").Value = 0
$ws.Range("R128").Value = ''

# Row 129: A 9282-2025 (was row 65)
$ws.Range("A129").Value = 'A 9282-2025'
$ws.Range("B129").Value = 45714.64293981482
$ws.Range("C129").Value = 46077
$ws.Range("D129").Value = 'VÄSTRA GÖTALANDS LÄN'
$ws.Range("E129").Value = 'STENUNGSUND'
$ws.Range("G129").Value = 1.5
$ws.Range("H129").Value = 0
$ws.Range("I129").Value = 0
$ws.Range("J129").Value = 0
$ws.Range("K129").Value = 0
$ws.Range("L129").Value = 0
$ws.Range("M129").Value = 0
$ws.Range("N129").Value = 0
$ws.Range("O129").Value = 0
$ws.Range("P129").Value = 0
$ws.Range("Q129").Value = 0
$ws.Range("R129").Value = ''

# Row 130: A 14981-2025 (was row 72)
$ws.Range("A130").Value = 'A 14981-2025'
$ws.Range("B130").Value = 45743.59576388889
$ws.Range("C130").Value = 46077
$ws.Range("D130").Value = 'VÄSTRA GÖTALANDS LÄN'
$ws.Range("E130").Value = 'STENUNGSUND'
$ws.Range("G130").Value = 4.8
$ws.Range("H130").Value = 0
$ws.Range("I130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("K130").Value = 0
$ws.Range("L130").Value = 0
$ws.Range("M130").Value = 0
$ws.Range("N130").Value = 0
$ws.Range("O130").Value = 0
$ws.Range("P130").Value = 0
$ws.Range("Q130").Value = 0
$ws.Range("R130").Value = ''
